# Auto-generated edit script: refresh scraped schedule data (Linea 141)
# Commit: "Horarios actualizados Linea 141 - 658"
#
# This mirrors a routine scrape refresh of the bus-arrival board:
#   - the "Ultima actualizacion" timestamp and "Total filas" counters on
#     each sheet are bumped to the new scrape time (17:38:19)
#   - twelve new arrivals captured at 17:38:19 are merged into sheet
#     "LP1912" (sorted by Hora_Llegada), pushing the row count from
#     367 to 379
#   - the matching subset of those new rows that belong to the "215"
#     family of lines is merged into sheet "LP1912-215" (39 -> 40)
#   - sheet "6203-6173" only gets its timestamp banner refreshed

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  (367 -> 379 data rows)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(6,1).Value = "03:45:25"
$ws1.Cells.Item(6,2).Value = "03:47"
$ws1.Cells.Item(6,3).Value = "14_ABASTO"
$ws1.Cells.Item(6,4).Value = 2
$ws1.Cells.Item(6,5).Value = "LP1912"
$ws1.Cells.Item(7,1).Value = "03:45:25"
$ws1.Cells.Item(7,2).Value = "04:01"
$ws1.Cells.Item(7,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(7,4).Value = 16
$ws1.Cells.Item(7,5).Value = "LP1912"
$ws1.Cells.Item(8,1).Value = "03:45:25"
$ws1.Cells.Item(8,2).Value = "04:46"
$ws1.Cells.Item(8,3).Value = "215A_EL PATO"
$ws1.Cells.Item(8,4).Value = 61
$ws1.Cells.Item(8,5).Value = "LP1912"
$ws1.Cells.Item(9,1).Value = "03:45:25"
$ws1.Cells.Item(9,2).Value = "04:53"
$ws1.Cells.Item(9,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(9,4).Value = 68
$ws1.Cells.Item(9,5).Value = "LP1912"
$ws1.Cells.Item(10,1).Value = "04:56:49"
$ws1.Cells.Item(10,2).Value = "05:13"
$ws1.Cells.Item(10,3).Value = "14_ABASTO"
$ws1.Cells.Item(10,4).Value = 17
$ws1.Cells.Item(10,5).Value = "LP1912"
$ws1.Cells.Item(11,1).Value = "03:45:25"
$ws1.Cells.Item(11,2).Value = "05:16"
$ws1.Cells.Item(11,3).Value = "17_ROMERO"
$ws1.Cells.Item(11,4).Value = 91
$ws1.Cells.Item(11,5).Value = "LP1912"
$ws1.Cells.Item(12,1).Value = "04:45:05"
$ws1.Cells.Item(12,2).Value = "05:16"
$ws1.Cells.Item(12,3).Value = "14_ABASTO"
$ws1.Cells.Item(12,4).Value = 31
$ws1.Cells.Item(12,5).Value = "LP1912"
$ws1.Cells.Item(13,1).Value = "03:45:25"
$ws1.Cells.Item(13,2).Value = "05:22"
$ws1.Cells.Item(13,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(13,4).Value = 97
$ws1.Cells.Item(13,5).Value = "LP1912"
$ws1.Cells.Item(14,1).Value = "05:26:08"
$ws1.Cells.Item(14,2).Value = "05:28"
$ws1.Cells.Item(14,3).Value = "14_ABASTO"
$ws1.Cells.Item(14,4).Value = 2
$ws1.Cells.Item(14,5).Value = "LP1912"
$ws1.Cells.Item(15,1).Value = "04:18:02"
$ws1.Cells.Item(15,2).Value = "05:34"
$ws1.Cells.Item(15,3).Value = "14_ABASTO"
$ws1.Cells.Item(15,4).Value = 76
$ws1.Cells.Item(15,5).Value = "LP1912"
$ws1.Cells.Item(16,1).Value = "03:45:25"
$ws1.Cells.Item(16,2).Value = "05:34"
$ws1.Cells.Item(16,3).Value = "215B_EL PATO"
$ws1.Cells.Item(16,4).Value = 109
$ws1.Cells.Item(16,5).Value = "LP1912"
$ws1.Cells.Item(17,1).Value = "04:18:02"
$ws1.Cells.Item(17,2).Value = "05:35"
$ws1.Cells.Item(17,3).Value = "215B_EL PATO"
$ws1.Cells.Item(17,4).Value = 77
$ws1.Cells.Item(17,5).Value = "LP1912"
$ws1.Cells.Item(18,1).Value = "03:45:25"
$ws1.Cells.Item(18,2).Value = "05:37"
$ws1.Cells.Item(18,3).Value = "14_ABASTO"
$ws1.Cells.Item(18,4).Value = 112
$ws1.Cells.Item(18,5).Value = "LP1912"
$ws1.Cells.Item(19,1).Value = "04:18:02"
$ws1.Cells.Item(19,2).Value = "05:46"
$ws1.Cells.Item(19,3).Value = "15_ABASTO"
$ws1.Cells.Item(19,4).Value = 88
$ws1.Cells.Item(19,5).Value = "LP1912"
$ws1.Cells.Item(20,1).Value = "04:45:05"
$ws1.Cells.Item(20,2).Value = "06:04"
$ws1.Cells.Item(20,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(20,4).Value = 79
$ws1.Cells.Item(20,5).Value = "LP1912"
$ws1.Cells.Item(21,1).Value = "04:18:02"
$ws1.Cells.Item(21,2).Value = "06:05"
$ws1.Cells.Item(21,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(21,4).Value = 107
$ws1.Cells.Item(21,5).Value = "LP1912"
$ws1.Cells.Item(22,1).Value = "04:56:49"
$ws1.Cells.Item(22,2).Value = "06:11"
$ws1.Cells.Item(22,3).Value = "215A_EL PATO"
$ws1.Cells.Item(22,4).Value = 75
$ws1.Cells.Item(22,5).Value = "LP1912"
$ws1.Cells.Item(23,1).Value = "04:18:02"
$ws1.Cells.Item(23,2).Value = "06:12"
$ws1.Cells.Item(23,3).Value = "215A_EL PATO"
$ws1.Cells.Item(23,4).Value = 114
$ws1.Cells.Item(23,5).Value = "LP1912"
$ws1.Cells.Item(24,1).Value = "04:18:02"
$ws1.Cells.Item(24,2).Value = "06:14"
$ws1.Cells.Item(24,3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(24,4).Value = 116
$ws1.Cells.Item(24,5).Value = "LP1912"
$ws1.Cells.Item(25,1).Value = "04:45:05"
$ws1.Cells.Item(25,2).Value = "06:21"
$ws1.Cells.Item(25,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(25,4).Value = 96
$ws1.Cells.Item(25,5).Value = "LP1912"
$ws1.Cells.Item(26,1).Value = "06:25:43"
$ws1.Cells.Item(26,2).Value = "06:26"
$ws1.Cells.Item(26,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(26,4).Value = 1
$ws1.Cells.Item(26,5).Value = "LP1912"
$ws1.Cells.Item(27,1).Value = "04:45:05"
$ws1.Cells.Item(27,2).Value = "06:27"
$ws1.Cells.Item(27,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(27,4).Value = 102
$ws1.Cells.Item(27,5).Value = "LP1912"
$ws1.Cells.Item(28,1).Value = "06:25:43"
$ws1.Cells.Item(28,2).Value = "06:28"
$ws1.Cells.Item(28,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(28,4).Value = 3
$ws1.Cells.Item(28,5).Value = "LP1912"
$ws1.Cells.Item(29,1).Value = "04:56:49"
$ws1.Cells.Item(29,2).Value = "06:29"
$ws1.Cells.Item(29,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(29,4).Value = 93
$ws1.Cells.Item(29,5).Value = "LP1912"
$ws1.Cells.Item(30,1).Value = "04:45:05"
$ws1.Cells.Item(30,2).Value = "06:30"
$ws1.Cells.Item(30,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(30,4).Value = 105
$ws1.Cells.Item(30,5).Value = "LP1912"
$ws1.Cells.Item(31,1).Value = "04:45:05"
$ws1.Cells.Item(31,2).Value = "06:31"
$ws1.Cells.Item(31,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(31,4).Value = 106
$ws1.Cells.Item(31,5).Value = "LP1912"
$ws1.Cells.Item(32,1).Value = "05:55:25"
$ws1.Cells.Item(32,2).Value = "06:44"
$ws1.Cells.Item(32,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(32,4).Value = 49
$ws1.Cells.Item(32,5).Value = "LP1912"
$ws1.Cells.Item(33,1).Value = "04:45:05"
$ws1.Cells.Item(33,2).Value = "06:44"
$ws1.Cells.Item(33,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(33,4).Value = 119
$ws1.Cells.Item(33,5).Value = "LP1912"
$ws1.Cells.Item(34,1).Value = "04:56:49"
$ws1.Cells.Item(34,2).Value = "06:46"
$ws1.Cells.Item(34,3).Value = "215C_EL PATO"
$ws1.Cells.Item(34,4).Value = 110
$ws1.Cells.Item(34,5).Value = "LP1912"
$ws1.Cells.Item(35,1).Value = "05:26:08"
$ws1.Cells.Item(35,2).Value = "06:47"
$ws1.Cells.Item(35,3).Value = "215C_EL PATO"
$ws1.Cells.Item(35,4).Value = 81
$ws1.Cells.Item(35,5).Value = "LP1912"
$ws1.Cells.Item(36,1).Value = "05:55:25"
$ws1.Cells.Item(36,2).Value = "06:59"
$ws1.Cells.Item(36,3).Value = "14_ABASTO"
$ws1.Cells.Item(36,4).Value = 64
$ws1.Cells.Item(36,5).Value = "LP1912"
$ws1.Cells.Item(37,1).Value = "05:26:08"
$ws1.Cells.Item(37,2).Value = "07:00"
$ws1.Cells.Item(37,3).Value = "14_ABASTO"
$ws1.Cells.Item(37,4).Value = 94
$ws1.Cells.Item(37,5).Value = "LP1912"
$ws1.Cells.Item(38,1).Value = "06:25:43"
$ws1.Cells.Item(38,2).Value = "07:01"
$ws1.Cells.Item(38,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(38,4).Value = 36
$ws1.Cells.Item(38,5).Value = "LP1912"
$ws1.Cells.Item(39,1).Value = "05:55:25"
$ws1.Cells.Item(39,2).Value = "07:04"
$ws1.Cells.Item(39,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(39,4).Value = 69
$ws1.Cells.Item(39,5).Value = "LP1912"
$ws1.Cells.Item(40,1).Value = "05:26:08"
$ws1.Cells.Item(40,2).Value = "07:05"
$ws1.Cells.Item(40,3).Value = "15_ABASTO"
$ws1.Cells.Item(40,4).Value = 99
$ws1.Cells.Item(40,5).Value = "LP1912"
$ws1.Cells.Item(41,1).Value = "05:26:08"
$ws1.Cells.Item(41,2).Value = "07:05"
$ws1.Cells.Item(41,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(41,4).Value = 99
$ws1.Cells.Item(41,5).Value = "LP1912"
$ws1.Cells.Item(42,1).Value = "05:26:08"
$ws1.Cells.Item(42,2).Value = "07:06"
$ws1.Cells.Item(42,3).Value = "10_OLMOS"
$ws1.Cells.Item(42,4).Value = 100
$ws1.Cells.Item(42,5).Value = "LP1912"
$ws1.Cells.Item(43,1).Value = "05:26:08"
$ws1.Cells.Item(43,2).Value = "07:07"
$ws1.Cells.Item(43,3).Value = "225_GOMEZ"
$ws1.Cells.Item(43,4).Value = 101
$ws1.Cells.Item(43,5).Value = "LP1912"
$ws1.Cells.Item(44,1).Value = "05:26:08"
$ws1.Cells.Item(44,2).Value = "07:11"
$ws1.Cells.Item(44,3).Value = "215A_EL PATO"
$ws1.Cells.Item(44,4).Value = 105
$ws1.Cells.Item(44,5).Value = "LP1912"
$ws1.Cells.Item(45,1).Value = "06:55:02"
$ws1.Cells.Item(45,2).Value = "07:12"
$ws1.Cells.Item(45,3).Value = "215A_EL PATO"
$ws1.Cells.Item(45,4).Value = 17
$ws1.Cells.Item(45,5).Value = "LP1912"
$ws1.Cells.Item(46,1).Value = "06:25:43"
$ws1.Cells.Item(46,2).Value = "07:14"
$ws1.Cells.Item(46,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(46,4).Value = 49
$ws1.Cells.Item(46,5).Value = "LP1912"
$ws1.Cells.Item(47,1).Value = "05:55:25"
$ws1.Cells.Item(47,2).Value = "07:15"
$ws1.Cells.Item(47,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(47,4).Value = 80
$ws1.Cells.Item(47,5).Value = "LP1912"
$ws1.Cells.Item(48,1).Value = "05:26:08"
$ws1.Cells.Item(48,2).Value = "07:16"
$ws1.Cells.Item(48,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(48,4).Value = 110
$ws1.Cells.Item(48,5).Value = "LP1912"
$ws1.Cells.Item(49,1).Value = "06:55:02"
$ws1.Cells.Item(49,2).Value = "07:17"
$ws1.Cells.Item(49,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(49,4).Value = 22
$ws1.Cells.Item(49,5).Value = "LP1912"
$ws1.Cells.Item(50,1).Value = "05:26:08"
$ws1.Cells.Item(50,2).Value = "07:21"
$ws1.Cells.Item(50,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(50,4).Value = 115
$ws1.Cells.Item(50,5).Value = "LP1912"
$ws1.Cells.Item(51,1).Value = "05:26:08"
$ws1.Cells.Item(51,2).Value = "07:23"
$ws1.Cells.Item(51,3).Value = "10_OLMOS"
$ws1.Cells.Item(51,4).Value = 117
$ws1.Cells.Item(51,5).Value = "LP1912"
$ws1.Cells.Item(52,1).Value = "05:55:25"
$ws1.Cells.Item(52,2).Value = "07:30"
$ws1.Cells.Item(52,3).Value = "10_OLMOS"
$ws1.Cells.Item(52,4).Value = 95
$ws1.Cells.Item(52,5).Value = "LP1912"
$ws1.Cells.Item(53,1).Value = "05:55:25"
$ws1.Cells.Item(53,2).Value = "07:31"
$ws1.Cells.Item(53,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(53,4).Value = 96
$ws1.Cells.Item(53,5).Value = "LP1912"
$ws1.Cells.Item(54,1).Value = "05:55:25"
$ws1.Cells.Item(54,2).Value = "07:31"
$ws1.Cells.Item(54,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(54,4).Value = 96
$ws1.Cells.Item(54,5).Value = "LP1912"
$ws1.Cells.Item(55,1).Value = "06:55:02"
$ws1.Cells.Item(55,2).Value = "07:32"
$ws1.Cells.Item(55,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(55,4).Value = 37
$ws1.Cells.Item(55,5).Value = "LP1912"
$ws1.Cells.Item(56,1).Value = "05:55:25"
$ws1.Cells.Item(56,2).Value = "07:32"
$ws1.Cells.Item(56,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(56,4).Value = 97
$ws1.Cells.Item(56,5).Value = "LP1912"
$ws1.Cells.Item(57,1).Value = "06:55:02"
$ws1.Cells.Item(57,2).Value = "07:32"
$ws1.Cells.Item(57,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(57,4).Value = 37
$ws1.Cells.Item(57,5).Value = "LP1912"
$ws1.Cells.Item(58,1).Value = "07:19:29"
$ws1.Cells.Item(58,2).Value = "07:35"
$ws1.Cells.Item(58,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(58,4).Value = 16
$ws1.Cells.Item(58,5).Value = "LP1912"
$ws1.Cells.Item(59,1).Value = "05:55:25"
$ws1.Cells.Item(59,2).Value = "07:36"
$ws1.Cells.Item(59,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(59,4).Value = 101
$ws1.Cells.Item(59,5).Value = "LP1912"
$ws1.Cells.Item(60,1).Value = "06:55:02"
$ws1.Cells.Item(60,2).Value = "07:37"
$ws1.Cells.Item(60,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(60,4).Value = 42
$ws1.Cells.Item(60,5).Value = "LP1912"
$ws1.Cells.Item(61,1).Value = "05:55:25"
$ws1.Cells.Item(61,2).Value = "07:39"
$ws1.Cells.Item(61,3).Value = "10_OLMOS"
$ws1.Cells.Item(61,4).Value = 104
$ws1.Cells.Item(61,5).Value = "LP1912"
$ws1.Cells.Item(62,1).Value = "07:19:29"
$ws1.Cells.Item(62,2).Value = "07:46"
$ws1.Cells.Item(62,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(62,4).Value = 27
$ws1.Cells.Item(62,5).Value = "LP1912"
$ws1.Cells.Item(63,1).Value = "05:55:25"
$ws1.Cells.Item(63,2).Value = "07:47"
$ws1.Cells.Item(63,3).Value = "14_ABASTO"
$ws1.Cells.Item(63,4).Value = 112
$ws1.Cells.Item(63,5).Value = "LP1912"
$ws1.Cells.Item(64,1).Value = "06:55:02"
$ws1.Cells.Item(64,2).Value = "07:48"
$ws1.Cells.Item(64,3).Value = "14_ABASTO"
$ws1.Cells.Item(64,4).Value = 53
$ws1.Cells.Item(64,5).Value = "LP1912"
$ws1.Cells.Item(65,1).Value = "07:50:16"
$ws1.Cells.Item(65,2).Value = "07:50"
$ws1.Cells.Item(65,3).Value = "10_OLMOS"
$ws1.Cells.Item(65,4).Value = 0
$ws1.Cells.Item(65,5).Value = "LP1912"
$ws1.Cells.Item(66,1).Value = "05:55:25"
$ws1.Cells.Item(66,2).Value = "07:51"
$ws1.Cells.Item(66,3).Value = "215D_EL PATO"
$ws1.Cells.Item(66,4).Value = 116
$ws1.Cells.Item(66,5).Value = "LP1912"
$ws1.Cells.Item(67,1).Value = "06:55:02"
$ws1.Cells.Item(67,2).Value = "07:52"
$ws1.Cells.Item(67,3).Value = "215D_EL PATO"
$ws1.Cells.Item(67,4).Value = 57
$ws1.Cells.Item(67,5).Value = "LP1912"
$ws1.Cells.Item(68,1).Value = "07:19:29"
$ws1.Cells.Item(68,2).Value = "07:59"
$ws1.Cells.Item(68,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(68,4).Value = 40
$ws1.Cells.Item(68,5).Value = "LP1912"
$ws1.Cells.Item(69,1).Value = "06:25:43"
$ws1.Cells.Item(69,2).Value = "08:01"
$ws1.Cells.Item(69,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(69,4).Value = 96
$ws1.Cells.Item(69,5).Value = "LP1912"
$ws1.Cells.Item(70,1).Value = "07:19:29"
$ws1.Cells.Item(70,2).Value = "08:03"
$ws1.Cells.Item(70,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(70,4).Value = 44
$ws1.Cells.Item(70,5).Value = "LP1912"
$ws1.Cells.Item(71,1).Value = "06:55:02"
$ws1.Cells.Item(71,2).Value = "08:03"
$ws1.Cells.Item(71,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(71,4).Value = 68
$ws1.Cells.Item(71,5).Value = "LP1912"
$ws1.Cells.Item(72,1).Value = "08:02:22"
$ws1.Cells.Item(72,2).Value = "08:05"
$ws1.Cells.Item(72,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(72,4).Value = 3
$ws1.Cells.Item(72,5).Value = "LP1912"
$ws1.Cells.Item(73,1).Value = "07:19:29"
$ws1.Cells.Item(73,2).Value = "08:10"
$ws1.Cells.Item(73,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(73,4).Value = 51
$ws1.Cells.Item(73,5).Value = "LP1912"
$ws1.Cells.Item(74,1).Value = "07:50:16"
$ws1.Cells.Item(74,2).Value = "08:11"
$ws1.Cells.Item(74,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(74,4).Value = 21
$ws1.Cells.Item(74,5).Value = "LP1912"
$ws1.Cells.Item(75,1).Value = "06:25:43"
$ws1.Cells.Item(75,2).Value = "08:12"
$ws1.Cells.Item(75,3).Value = "15_ABASTO"
$ws1.Cells.Item(75,4).Value = 107
$ws1.Cells.Item(75,5).Value = "LP1912"
$ws1.Cells.Item(76,1).Value = "07:50:16"
$ws1.Cells.Item(76,2).Value = "08:13"
$ws1.Cells.Item(76,3).Value = "10_OLMOS"
$ws1.Cells.Item(76,4).Value = 23
$ws1.Cells.Item(76,5).Value = "LP1912"
$ws1.Cells.Item(77,1).Value = "06:55:02"
$ws1.Cells.Item(77,2).Value = "08:21"
$ws1.Cells.Item(77,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(77,4).Value = 86
$ws1.Cells.Item(77,5).Value = "LP1912"
$ws1.Cells.Item(78,1).Value = "06:25:43"
$ws1.Cells.Item(78,2).Value = "08:22"
$ws1.Cells.Item(78,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(78,4).Value = 117
$ws1.Cells.Item(78,5).Value = "LP1912"
$ws1.Cells.Item(79,1).Value = "06:55:02"
$ws1.Cells.Item(79,2).Value = "08:23"
$ws1.Cells.Item(79,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(79,4).Value = 88
$ws1.Cells.Item(79,5).Value = "LP1912"
$ws1.Cells.Item(80,1).Value = "06:25:43"
$ws1.Cells.Item(80,2).Value = "08:23"
$ws1.Cells.Item(80,3).Value = "215B_EL PATO"
$ws1.Cells.Item(80,4).Value = 118
$ws1.Cells.Item(80,5).Value = "LP1912"
$ws1.Cells.Item(81,1).Value = "06:55:02"
$ws1.Cells.Item(81,2).Value = "08:27"
$ws1.Cells.Item(81,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(81,4).Value = 92
$ws1.Cells.Item(81,5).Value = "LP1912"
$ws1.Cells.Item(82,1).Value = "07:50:16"
$ws1.Cells.Item(82,2).Value = "08:30"
$ws1.Cells.Item(82,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(82,4).Value = 40
$ws1.Cells.Item(82,5).Value = "LP1912"
$ws1.Cells.Item(83,1).Value = "08:02:22"
$ws1.Cells.Item(83,2).Value = "08:33"
$ws1.Cells.Item(83,3).Value = "10_OLMOS"
$ws1.Cells.Item(83,4).Value = 31
$ws1.Cells.Item(83,5).Value = "LP1912"
$ws1.Cells.Item(84,1).Value = "08:02:22"
$ws1.Cells.Item(84,2).Value = "08:34"
$ws1.Cells.Item(84,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(84,4).Value = 32
$ws1.Cells.Item(84,5).Value = "LP1912"
$ws1.Cells.Item(85,1).Value = "08:32:09"
$ws1.Cells.Item(85,2).Value = "08:37"
$ws1.Cells.Item(85,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(85,4).Value = 5
$ws1.Cells.Item(85,5).Value = "LP1912"
$ws1.Cells.Item(86,1).Value = "06:55:02"
$ws1.Cells.Item(86,2).Value = "08:42"
$ws1.Cells.Item(86,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(86,4).Value = 107
$ws1.Cells.Item(86,5).Value = "LP1912"
$ws1.Cells.Item(87,1).Value = "07:19:29"
$ws1.Cells.Item(87,2).Value = "08:43"
$ws1.Cells.Item(87,3).Value = "14_ABASTO"
$ws1.Cells.Item(87,4).Value = 84
$ws1.Cells.Item(87,5).Value = "LP1912"
$ws1.Cells.Item(88,1).Value = "07:50:16"
$ws1.Cells.Item(88,2).Value = "08:44"
$ws1.Cells.Item(88,3).Value = "14_ABASTO"
$ws1.Cells.Item(88,4).Value = 54
$ws1.Cells.Item(88,5).Value = "LP1912"
$ws1.Cells.Item(89,1).Value = "08:32:09"
$ws1.Cells.Item(89,2).Value = "08:53"
$ws1.Cells.Item(89,3).Value = "10_OLMOS"
$ws1.Cells.Item(89,4).Value = 21
$ws1.Cells.Item(89,5).Value = "LP1912"
$ws1.Cells.Item(90,1).Value = "06:55:02"
$ws1.Cells.Item(90,2).Value = "08:54"
$ws1.Cells.Item(90,3).Value = "17_ROMERO"
$ws1.Cells.Item(90,4).Value = 119
$ws1.Cells.Item(90,5).Value = "LP1912"
$ws1.Cells.Item(91,1).Value = "07:19:29"
$ws1.Cells.Item(91,2).Value = "09:01"
$ws1.Cells.Item(91,3).Value = "215A_EL PATO"
$ws1.Cells.Item(91,4).Value = 102
$ws1.Cells.Item(91,5).Value = "LP1912"
$ws1.Cells.Item(92,1).Value = "07:50:16"
$ws1.Cells.Item(92,2).Value = "09:02"
$ws1.Cells.Item(92,3).Value = "215A_EL PATO"
$ws1.Cells.Item(92,4).Value = 72
$ws1.Cells.Item(92,5).Value = "LP1912"
$ws1.Cells.Item(93,1).Value = "08:02:22"
$ws1.Cells.Item(93,2).Value = "09:03"
$ws1.Cells.Item(93,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(93,4).Value = 61
$ws1.Cells.Item(93,5).Value = "LP1912"
$ws1.Cells.Item(94,1).Value = "08:32:09"
$ws1.Cells.Item(94,2).Value = "09:04"
$ws1.Cells.Item(94,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(94,4).Value = 32
$ws1.Cells.Item(94,5).Value = "LP1912"
$ws1.Cells.Item(95,1).Value = "08:32:09"
$ws1.Cells.Item(95,2).Value = "09:05"
$ws1.Cells.Item(95,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(95,4).Value = 33
$ws1.Cells.Item(95,5).Value = "LP1912"
$ws1.Cells.Item(96,1).Value = "07:19:29"
$ws1.Cells.Item(96,2).Value = "09:10"
$ws1.Cells.Item(96,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(96,4).Value = 111
$ws1.Cells.Item(96,5).Value = "LP1912"
$ws1.Cells.Item(97,1).Value = "07:50:16"
$ws1.Cells.Item(97,2).Value = "09:11"
$ws1.Cells.Item(97,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(97,4).Value = 81
$ws1.Cells.Item(97,5).Value = "LP1912"
$ws1.Cells.Item(98,1).Value = "08:32:09"
$ws1.Cells.Item(98,2).Value = "09:13"
$ws1.Cells.Item(98,3).Value = "10_OLMOS"
$ws1.Cells.Item(98,4).Value = 41
$ws1.Cells.Item(98,5).Value = "LP1912"
$ws1.Cells.Item(99,1).Value = "07:19:29"
$ws1.Cells.Item(99,2).Value = "09:16"
$ws1.Cells.Item(99,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(99,4).Value = 117
$ws1.Cells.Item(99,5).Value = "LP1912"
$ws1.Cells.Item(100,1).Value = "07:50:16"
$ws1.Cells.Item(100,2).Value = "09:17"
$ws1.Cells.Item(100,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(100,4).Value = 87
$ws1.Cells.Item(100,5).Value = "LP1912"
$ws1.Cells.Item(101,1).Value = "07:50:16"
$ws1.Cells.Item(101,2).Value = "09:21"
$ws1.Cells.Item(101,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(101,4).Value = 91
$ws1.Cells.Item(101,5).Value = "LP1912"
$ws1.Cells.Item(102,1).Value = "08:02:22"
$ws1.Cells.Item(102,2).Value = "09:22"
$ws1.Cells.Item(102,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(102,4).Value = 80
$ws1.Cells.Item(102,5).Value = "LP1912"
$ws1.Cells.Item(103,1).Value = "08:32:09"
$ws1.Cells.Item(103,2).Value = "09:23"
$ws1.Cells.Item(103,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(103,4).Value = 51
$ws1.Cells.Item(103,5).Value = "LP1912"
$ws1.Cells.Item(104,1).Value = "07:50:16"
$ws1.Cells.Item(104,2).Value = "09:23"
$ws1.Cells.Item(104,3).Value = "17_ROMERO"
$ws1.Cells.Item(104,4).Value = 93
$ws1.Cells.Item(104,5).Value = "LP1912"
$ws1.Cells.Item(105,1).Value = "08:02:22"
$ws1.Cells.Item(105,2).Value = "09:23"
$ws1.Cells.Item(105,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(105,4).Value = 81
$ws1.Cells.Item(105,5).Value = "LP1912"
$ws1.Cells.Item(106,1).Value = "07:50:16"
$ws1.Cells.Item(106,2).Value = "09:24"
$ws1.Cells.Item(106,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(106,4).Value = 94
$ws1.Cells.Item(106,5).Value = "LP1912"
$ws1.Cells.Item(107,1).Value = "07:50:16"
$ws1.Cells.Item(107,2).Value = "09:28"
$ws1.Cells.Item(107,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(107,4).Value = 98
$ws1.Cells.Item(107,5).Value = "LP1912"
$ws1.Cells.Item(108,1).Value = "07:50:16"
$ws1.Cells.Item(108,2).Value = "09:32"
$ws1.Cells.Item(108,3).Value = "15_ABASTO"
$ws1.Cells.Item(108,4).Value = 102
$ws1.Cells.Item(108,5).Value = "LP1912"
$ws1.Cells.Item(109,1).Value = "07:50:16"
$ws1.Cells.Item(109,2).Value = "09:33"
$ws1.Cells.Item(109,3).Value = "10_OLMOS"
$ws1.Cells.Item(109,4).Value = 103
$ws1.Cells.Item(109,5).Value = "LP1912"
$ws1.Cells.Item(110,1).Value = "08:56:29"
$ws1.Cells.Item(110,2).Value = "09:34"
$ws1.Cells.Item(110,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(110,4).Value = 38
$ws1.Cells.Item(110,5).Value = "LP1912"
$ws1.Cells.Item(111,1).Value = "08:56:29"
$ws1.Cells.Item(111,2).Value = "09:34"
$ws1.Cells.Item(111,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(111,4).Value = 38
$ws1.Cells.Item(111,5).Value = "LP1912"
$ws1.Cells.Item(112,1).Value = "08:32:09"
$ws1.Cells.Item(112,2).Value = "09:35"
$ws1.Cells.Item(112,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(112,4).Value = 63
$ws1.Cells.Item(112,5).Value = "LP1912"
$ws1.Cells.Item(113,1).Value = "08:48:08"
$ws1.Cells.Item(113,2).Value = "09:35"
$ws1.Cells.Item(113,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(113,4).Value = 47
$ws1.Cells.Item(113,5).Value = "LP1912"
$ws1.Cells.Item(114,1).Value = "09:35:26"
$ws1.Cells.Item(114,2).Value = "09:39"
$ws1.Cells.Item(114,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(114,4).Value = 4
$ws1.Cells.Item(114,5).Value = "LP1912"
$ws1.Cells.Item(115,1).Value = "07:50:16"
$ws1.Cells.Item(115,2).Value = "09:42"
$ws1.Cells.Item(115,3).Value = "215C_EL PATO"
$ws1.Cells.Item(115,4).Value = 112
$ws1.Cells.Item(115,5).Value = "LP1912"
$ws1.Cells.Item(116,1).Value = "08:02:22"
$ws1.Cells.Item(116,2).Value = "09:43"
$ws1.Cells.Item(116,3).Value = "14_ABASTO"
$ws1.Cells.Item(116,4).Value = 101
$ws1.Cells.Item(116,5).Value = "LP1912"
$ws1.Cells.Item(117,1).Value = "07:50:16"
$ws1.Cells.Item(117,2).Value = "09:44"
$ws1.Cells.Item(117,3).Value = "14_ABASTO"
$ws1.Cells.Item(117,4).Value = 114
$ws1.Cells.Item(117,5).Value = "LP1912"
$ws1.Cells.Item(118,1).Value = "09:35:26"
$ws1.Cells.Item(118,2).Value = "09:46"
$ws1.Cells.Item(118,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(118,4).Value = 11
$ws1.Cells.Item(118,5).Value = "LP1912"
$ws1.Cells.Item(119,1).Value = "08:32:09"
$ws1.Cells.Item(119,2).Value = "09:52"
$ws1.Cells.Item(119,3).Value = "15_ABASTO"
$ws1.Cells.Item(119,4).Value = 80
$ws1.Cells.Item(119,5).Value = "LP1912"
$ws1.Cells.Item(120,1).Value = "08:56:29"
$ws1.Cells.Item(120,2).Value = "09:53"
$ws1.Cells.Item(120,3).Value = "10_OLMOS"
$ws1.Cells.Item(120,4).Value = 57
$ws1.Cells.Item(120,5).Value = "LP1912"
$ws1.Cells.Item(121,1).Value = "09:35:26"
$ws1.Cells.Item(121,2).Value = "09:58"
$ws1.Cells.Item(121,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(121,4).Value = 23
$ws1.Cells.Item(121,5).Value = "LP1912"
$ws1.Cells.Item(122,1).Value = "09:35:26"
$ws1.Cells.Item(122,2).Value = "10:03"
$ws1.Cells.Item(122,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(122,4).Value = 28
$ws1.Cells.Item(122,5).Value = "LP1912"
$ws1.Cells.Item(123,1).Value = "08:56:29"
$ws1.Cells.Item(123,2).Value = "10:10"
$ws1.Cells.Item(123,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(123,4).Value = 74
$ws1.Cells.Item(123,5).Value = "LP1912"
$ws1.Cells.Item(124,1).Value = "08:32:09"
$ws1.Cells.Item(124,2).Value = "10:11"
$ws1.Cells.Item(124,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(124,4).Value = 99
$ws1.Cells.Item(124,5).Value = "LP1912"
$ws1.Cells.Item(125,1).Value = "09:35:26"
$ws1.Cells.Item(125,2).Value = "10:12"
$ws1.Cells.Item(125,3).Value = "15_ABASTO"
$ws1.Cells.Item(125,4).Value = 37
$ws1.Cells.Item(125,5).Value = "LP1912"
$ws1.Cells.Item(126,1).Value = "09:35:26"
$ws1.Cells.Item(126,2).Value = "10:13"
$ws1.Cells.Item(126,3).Value = "10_OLMOS"
$ws1.Cells.Item(126,4).Value = 38
$ws1.Cells.Item(126,5).Value = "LP1912"
$ws1.Cells.Item(127,1).Value = "08:32:09"
$ws1.Cells.Item(127,2).Value = "10:21"
$ws1.Cells.Item(127,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(127,4).Value = 109
$ws1.Cells.Item(127,5).Value = "LP1912"
$ws1.Cells.Item(128,1).Value = "08:32:09"
$ws1.Cells.Item(128,2).Value = "10:22"
$ws1.Cells.Item(128,3).Value = "17_ROMERO"
$ws1.Cells.Item(128,4).Value = 110
$ws1.Cells.Item(128,5).Value = "LP1912"
$ws1.Cells.Item(129,1).Value = "09:35:26"
$ws1.Cells.Item(129,2).Value = "10:23"
$ws1.Cells.Item(129,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(129,4).Value = 48
$ws1.Cells.Item(129,5).Value = "LP1912"
$ws1.Cells.Item(130,1).Value = "08:56:29"
$ws1.Cells.Item(130,2).Value = "10:26"
$ws1.Cells.Item(130,3).Value = "215A_EL PATO"
$ws1.Cells.Item(130,4).Value = 90
$ws1.Cells.Item(130,5).Value = "LP1912"
$ws1.Cells.Item(131,1).Value = "08:32:09"
$ws1.Cells.Item(131,2).Value = "10:27"
$ws1.Cells.Item(131,3).Value = "215A_EL PATO"
$ws1.Cells.Item(131,4).Value = 115
$ws1.Cells.Item(131,5).Value = "LP1912"
$ws1.Cells.Item(132,1).Value = "10:29:57"
$ws1.Cells.Item(132,2).Value = "10:29"
$ws1.Cells.Item(132,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(132,4).Value = 0
$ws1.Cells.Item(132,5).Value = "LP1912"
$ws1.Cells.Item(133,1).Value = "10:29:57"
$ws1.Cells.Item(133,2).Value = "10:31"
$ws1.Cells.Item(133,3).Value = "10_OLMOS"
$ws1.Cells.Item(133,4).Value = 2
$ws1.Cells.Item(133,5).Value = "LP1912"
$ws1.Cells.Item(134,1).Value = "10:29:57"
$ws1.Cells.Item(134,2).Value = "10:34"
$ws1.Cells.Item(134,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(134,4).Value = 5
$ws1.Cells.Item(134,5).Value = "LP1912"
$ws1.Cells.Item(135,1).Value = "09:35:26"
$ws1.Cells.Item(135,2).Value = "10:34"
$ws1.Cells.Item(135,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(135,4).Value = 59
$ws1.Cells.Item(135,5).Value = "LP1912"
$ws1.Cells.Item(136,1).Value = "10:29:57"
$ws1.Cells.Item(136,2).Value = "10:39"
$ws1.Cells.Item(136,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(136,4).Value = 10
$ws1.Cells.Item(136,5).Value = "LP1912"
$ws1.Cells.Item(137,1).Value = "10:29:57"
$ws1.Cells.Item(137,2).Value = "10:41"
$ws1.Cells.Item(137,3).Value = "17_ROMERO"
$ws1.Cells.Item(137,4).Value = 12
$ws1.Cells.Item(137,5).Value = "LP1912"
$ws1.Cells.Item(138,1).Value = "08:48:08"
$ws1.Cells.Item(138,2).Value = "10:42"
$ws1.Cells.Item(138,3).Value = "17_ROMERO"
$ws1.Cells.Item(138,4).Value = 114
$ws1.Cells.Item(138,5).Value = "LP1912"
$ws1.Cells.Item(139,1).Value = "08:56:29"
$ws1.Cells.Item(139,2).Value = "10:43"
$ws1.Cells.Item(139,3).Value = "14_ABASTO"
$ws1.Cells.Item(139,4).Value = 107
$ws1.Cells.Item(139,5).Value = "LP1912"
$ws1.Cells.Item(140,1).Value = "08:48:08"
$ws1.Cells.Item(140,2).Value = "10:44"
$ws1.Cells.Item(140,3).Value = "14_ABASTO"
$ws1.Cells.Item(140,4).Value = 116
$ws1.Cells.Item(140,5).Value = "LP1912"
$ws1.Cells.Item(141,1).Value = "10:29:57"
$ws1.Cells.Item(141,2).Value = "10:51"
$ws1.Cells.Item(141,3).Value = "15_ABASTO"
$ws1.Cells.Item(141,4).Value = 22
$ws1.Cells.Item(141,5).Value = "LP1912"
$ws1.Cells.Item(142,1).Value = "10:29:57"
$ws1.Cells.Item(142,2).Value = "10:52"
$ws1.Cells.Item(142,3).Value = "10_OLMOS"
$ws1.Cells.Item(142,4).Value = 23
$ws1.Cells.Item(142,5).Value = "LP1912"
$ws1.Cells.Item(143,1).Value = "09:35:26"
$ws1.Cells.Item(143,2).Value = "10:54"
$ws1.Cells.Item(143,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(143,4).Value = 79
$ws1.Cells.Item(143,5).Value = "LP1912"
$ws1.Cells.Item(144,1).Value = "10:29:57"
$ws1.Cells.Item(144,2).Value = "10:56"
$ws1.Cells.Item(144,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(144,4).Value = 27
$ws1.Cells.Item(144,5).Value = "LP1912"
$ws1.Cells.Item(145,1).Value = "10:59:49"
$ws1.Cells.Item(145,2).Value = "10:59"
$ws1.Cells.Item(145,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(145,4).Value = 0
$ws1.Cells.Item(145,5).Value = "LP1912"
$ws1.Cells.Item(146,1).Value = "10:29:57"
$ws1.Cells.Item(146,2).Value = "11:01"
$ws1.Cells.Item(146,3).Value = "215C_EL PATO"
$ws1.Cells.Item(146,4).Value = 32
$ws1.Cells.Item(146,5).Value = "LP1912"
$ws1.Cells.Item(147,1).Value = "09:35:26"
$ws1.Cells.Item(147,2).Value = "11:02"
$ws1.Cells.Item(147,3).Value = "215C_EL PATO"
$ws1.Cells.Item(147,4).Value = 87
$ws1.Cells.Item(147,5).Value = "LP1912"
$ws1.Cells.Item(148,1).Value = "10:29:57"
$ws1.Cells.Item(148,2).Value = "11:03"
$ws1.Cells.Item(148,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(148,4).Value = 34
$ws1.Cells.Item(148,5).Value = "LP1912"
$ws1.Cells.Item(149,1).Value = "10:29:57"
$ws1.Cells.Item(149,2).Value = "11:04"
$ws1.Cells.Item(149,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(149,4).Value = 35
$ws1.Cells.Item(149,5).Value = "LP1912"
$ws1.Cells.Item(150,1).Value = "10:59:49"
$ws1.Cells.Item(150,2).Value = "11:06"
$ws1.Cells.Item(150,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(150,4).Value = 7
$ws1.Cells.Item(150,5).Value = "LP1912"
$ws1.Cells.Item(151,1).Value = "09:35:26"
$ws1.Cells.Item(151,2).Value = "11:06"
$ws1.Cells.Item(151,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(151,4).Value = 91
$ws1.Cells.Item(151,5).Value = "LP1912"
$ws1.Cells.Item(152,1).Value = "10:29:57"
$ws1.Cells.Item(152,2).Value = "11:11"
$ws1.Cells.Item(152,3).Value = "15_ABASTO"
$ws1.Cells.Item(152,4).Value = 42
$ws1.Cells.Item(152,5).Value = "LP1912"
$ws1.Cells.Item(153,1).Value = "10:59:49"
$ws1.Cells.Item(153,2).Value = "11:11"
$ws1.Cells.Item(153,3).Value = "10_OLMOS"
$ws1.Cells.Item(153,4).Value = 12
$ws1.Cells.Item(153,5).Value = "LP1912"
$ws1.Cells.Item(154,1).Value = "10:59:49"
$ws1.Cells.Item(154,2).Value = "11:12"
$ws1.Cells.Item(154,3).Value = "15_ABASTO"
$ws1.Cells.Item(154,4).Value = 13
$ws1.Cells.Item(154,5).Value = "LP1912"
$ws1.Cells.Item(155,1).Value = "09:35:26"
$ws1.Cells.Item(155,2).Value = "11:19"
$ws1.Cells.Item(155,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(155,4).Value = 104
$ws1.Cells.Item(155,5).Value = "LP1912"
$ws1.Cells.Item(156,1).Value = "10:29:57"
$ws1.Cells.Item(156,2).Value = "11:20"
$ws1.Cells.Item(156,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(156,4).Value = 51
$ws1.Cells.Item(156,5).Value = "LP1912"
$ws1.Cells.Item(157,1).Value = "09:35:26"
$ws1.Cells.Item(157,2).Value = "11:21"
$ws1.Cells.Item(157,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(157,4).Value = 106
$ws1.Cells.Item(157,5).Value = "LP1912"
$ws1.Cells.Item(158,1).Value = "10:29:57"
$ws1.Cells.Item(158,2).Value = "11:26"
$ws1.Cells.Item(158,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(158,4).Value = 57
$ws1.Cells.Item(158,5).Value = "LP1912"
$ws1.Cells.Item(159,1).Value = "09:35:26"
$ws1.Cells.Item(159,2).Value = "11:27"
$ws1.Cells.Item(159,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(159,4).Value = 112
$ws1.Cells.Item(159,5).Value = "LP1912"
$ws1.Cells.Item(160,1).Value = "11:30:45"
$ws1.Cells.Item(160,2).Value = "11:30"
$ws1.Cells.Item(160,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(160,4).Value = 0
$ws1.Cells.Item(160,5).Value = "LP1912"
$ws1.Cells.Item(161,1).Value = "11:30:45"
$ws1.Cells.Item(161,2).Value = "11:31"
$ws1.Cells.Item(161,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(161,4).Value = 1
$ws1.Cells.Item(161,5).Value = "LP1912"
$ws1.Cells.Item(162,1).Value = "10:29:57"
$ws1.Cells.Item(162,2).Value = "11:31"
$ws1.Cells.Item(162,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(162,4).Value = 62
$ws1.Cells.Item(162,5).Value = "LP1912"
$ws1.Cells.Item(163,1).Value = "09:35:26"
$ws1.Cells.Item(163,2).Value = "11:32"
$ws1.Cells.Item(163,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(163,4).Value = 117
$ws1.Cells.Item(163,5).Value = "LP1912"
$ws1.Cells.Item(164,1).Value = "10:59:49"
$ws1.Cells.Item(164,2).Value = "11:34"
$ws1.Cells.Item(164,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(164,4).Value = 35
$ws1.Cells.Item(164,5).Value = "LP1912"
$ws1.Cells.Item(165,1).Value = "10:29:57"
$ws1.Cells.Item(165,2).Value = "11:35"
$ws1.Cells.Item(165,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(165,4).Value = 66
$ws1.Cells.Item(165,5).Value = "LP1912"
$ws1.Cells.Item(166,1).Value = "10:29:57"
$ws1.Cells.Item(166,2).Value = "11:40"
$ws1.Cells.Item(166,3).Value = "10_OLMOS"
$ws1.Cells.Item(166,4).Value = 71
$ws1.Cells.Item(166,5).Value = "LP1912"
$ws1.Cells.Item(167,1).Value = "10:29:57"
$ws1.Cells.Item(167,2).Value = "11:41"
$ws1.Cells.Item(167,3).Value = "17_ROMERO"
$ws1.Cells.Item(167,4).Value = 72
$ws1.Cells.Item(167,5).Value = "LP1912"
$ws1.Cells.Item(168,1).Value = "10:59:49"
$ws1.Cells.Item(168,2).Value = "11:42"
$ws1.Cells.Item(168,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(168,4).Value = 43
$ws1.Cells.Item(168,5).Value = "LP1912"
$ws1.Cells.Item(169,1).Value = "10:59:49"
$ws1.Cells.Item(169,2).Value = "11:43"
$ws1.Cells.Item(169,3).Value = "10_OLMOS"
$ws1.Cells.Item(169,4).Value = 44
$ws1.Cells.Item(169,5).Value = "LP1912"
$ws1.Cells.Item(170,1).Value = "11:30:45"
$ws1.Cells.Item(170,2).Value = "11:44"
$ws1.Cells.Item(170,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(170,4).Value = 14
$ws1.Cells.Item(170,5).Value = "LP1912"
$ws1.Cells.Item(171,1).Value = "10:29:57"
$ws1.Cells.Item(171,2).Value = "11:50"
$ws1.Cells.Item(171,3).Value = "215B_EL PATO"
$ws1.Cells.Item(171,4).Value = 81
$ws1.Cells.Item(171,5).Value = "LP1912"
$ws1.Cells.Item(172,1).Value = "10:59:49"
$ws1.Cells.Item(172,2).Value = "11:51"
$ws1.Cells.Item(172,3).Value = "215B_EL PATO"
$ws1.Cells.Item(172,4).Value = 52
$ws1.Cells.Item(172,5).Value = "LP1912"
$ws1.Cells.Item(173,1).Value = "10:59:49"
$ws1.Cells.Item(173,2).Value = "11:52"
$ws1.Cells.Item(173,3).Value = "15_ABASTO"
$ws1.Cells.Item(173,4).Value = 53
$ws1.Cells.Item(173,5).Value = "LP1912"
$ws1.Cells.Item(174,1).Value = "11:56:55"
$ws1.Cells.Item(174,2).Value = "11:56"
$ws1.Cells.Item(174,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(174,4).Value = 0
$ws1.Cells.Item(174,5).Value = "LP1912"
$ws1.Cells.Item(175,1).Value = "10:29:57"
$ws1.Cells.Item(175,2).Value = "11:58"
$ws1.Cells.Item(175,3).Value = "225_GOMEZ"
$ws1.Cells.Item(175,4).Value = 89
$ws1.Cells.Item(175,5).Value = "LP1912"
$ws1.Cells.Item(176,1).Value = "11:30:45"
$ws1.Cells.Item(176,2).Value = "11:59"
$ws1.Cells.Item(176,3).Value = "225_GOMEZ"
$ws1.Cells.Item(176,4).Value = 29
$ws1.Cells.Item(176,5).Value = "LP1912"
$ws1.Cells.Item(177,1).Value = "10:29:57"
$ws1.Cells.Item(177,2).Value = "12:01"
$ws1.Cells.Item(177,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(177,4).Value = 92
$ws1.Cells.Item(177,5).Value = "LP1912"
$ws1.Cells.Item(178,1).Value = "10:59:49"
$ws1.Cells.Item(178,2).Value = "12:02"
$ws1.Cells.Item(178,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(178,4).Value = 63
$ws1.Cells.Item(178,5).Value = "LP1912"
$ws1.Cells.Item(179,1).Value = "11:30:45"
$ws1.Cells.Item(179,2).Value = "12:04"
$ws1.Cells.Item(179,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(179,4).Value = 34
$ws1.Cells.Item(179,5).Value = "LP1912"
$ws1.Cells.Item(180,1).Value = "10:29:57"
$ws1.Cells.Item(180,2).Value = "12:06"
$ws1.Cells.Item(180,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(180,4).Value = 97
$ws1.Cells.Item(180,5).Value = "LP1912"
$ws1.Cells.Item(181,1).Value = "10:59:49"
$ws1.Cells.Item(181,2).Value = "12:06"
$ws1.Cells.Item(181,3).Value = "14_ABASTO"
$ws1.Cells.Item(181,4).Value = 67
$ws1.Cells.Item(181,5).Value = "LP1912"
$ws1.Cells.Item(182,1).Value = "10:59:49"
$ws1.Cells.Item(182,2).Value = "12:10"
$ws1.Cells.Item(182,3).Value = "10_OLMOS"
$ws1.Cells.Item(182,4).Value = 71
$ws1.Cells.Item(182,5).Value = "LP1912"
$ws1.Cells.Item(183,1).Value = "11:56:55"
$ws1.Cells.Item(183,2).Value = "12:12"
$ws1.Cells.Item(183,3).Value = "10_OLMOS"
$ws1.Cells.Item(183,4).Value = 16
$ws1.Cells.Item(183,5).Value = "LP1912"
$ws1.Cells.Item(184,1).Value = "10:29:57"
$ws1.Cells.Item(184,2).Value = "12:13"
$ws1.Cells.Item(184,3).Value = "17_ROMERO"
$ws1.Cells.Item(184,4).Value = 104
$ws1.Cells.Item(184,5).Value = "LP1912"
$ws1.Cells.Item(185,1).Value = "11:56:55"
$ws1.Cells.Item(185,2).Value = "12:14"
$ws1.Cells.Item(185,3).Value = "17_ROMERO"
$ws1.Cells.Item(185,4).Value = 18
$ws1.Cells.Item(185,5).Value = "LP1912"
$ws1.Cells.Item(186,1).Value = "10:29:57"
$ws1.Cells.Item(186,2).Value = "12:15"
$ws1.Cells.Item(186,3).Value = "14_ABASTO"
$ws1.Cells.Item(186,4).Value = 106
$ws1.Cells.Item(186,5).Value = "LP1912"
$ws1.Cells.Item(187,1).Value = "10:29:57"
$ws1.Cells.Item(187,2).Value = "12:20"
$ws1.Cells.Item(187,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(187,4).Value = 111
$ws1.Cells.Item(187,5).Value = "LP1912"
$ws1.Cells.Item(188,1).Value = "10:59:49"
$ws1.Cells.Item(188,2).Value = "12:20"
$ws1.Cells.Item(188,3).Value = "14_ABASTO"
$ws1.Cells.Item(188,4).Value = 81
$ws1.Cells.Item(188,5).Value = "LP1912"
$ws1.Cells.Item(189,1).Value = "10:29:57"
$ws1.Cells.Item(189,2).Value = "12:20"
$ws1.Cells.Item(189,3).Value = "215A_EL PATO"
$ws1.Cells.Item(189,4).Value = 111
$ws1.Cells.Item(189,5).Value = "LP1912"
$ws1.Cells.Item(190,1).Value = "10:59:49"
$ws1.Cells.Item(190,2).Value = "12:21"
$ws1.Cells.Item(190,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(190,4).Value = 82
$ws1.Cells.Item(190,5).Value = "LP1912"
$ws1.Cells.Item(191,1).Value = "12:21:08"
$ws1.Cells.Item(191,2).Value = "12:21"
$ws1.Cells.Item(191,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(191,4).Value = 0
$ws1.Cells.Item(191,5).Value = "LP1912"
$ws1.Cells.Item(192,1).Value = "12:21:08"
$ws1.Cells.Item(192,2).Value = "12:21"
$ws1.Cells.Item(192,3).Value = "215A_EL PATO"
$ws1.Cells.Item(192,4).Value = 0
$ws1.Cells.Item(192,5).Value = "LP1912"
$ws1.Cells.Item(193,1).Value = "10:59:49"
$ws1.Cells.Item(193,2).Value = "12:30"
$ws1.Cells.Item(193,3).Value = "17_ROMERO"
$ws1.Cells.Item(193,4).Value = 91
$ws1.Cells.Item(193,5).Value = "LP1912"
$ws1.Cells.Item(194,1).Value = "11:56:55"
$ws1.Cells.Item(194,2).Value = "12:34"
$ws1.Cells.Item(194,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(194,4).Value = 38
$ws1.Cells.Item(194,5).Value = "LP1912"
$ws1.Cells.Item(195,1).Value = "11:56:55"
$ws1.Cells.Item(195,2).Value = "12:34"
$ws1.Cells.Item(195,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(195,4).Value = 38
$ws1.Cells.Item(195,5).Value = "LP1912"
$ws1.Cells.Item(196,1).Value = "12:21:08"
$ws1.Cells.Item(196,2).Value = "12:35"
$ws1.Cells.Item(196,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(196,4).Value = 14
$ws1.Cells.Item(196,5).Value = "LP1912"
$ws1.Cells.Item(197,1).Value = "12:21:08"
$ws1.Cells.Item(197,2).Value = "12:35"
$ws1.Cells.Item(197,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(197,4).Value = 14
$ws1.Cells.Item(197,5).Value = "LP1912"
$ws1.Cells.Item(198,1).Value = "10:59:49"
$ws1.Cells.Item(198,2).Value = "12:36"
$ws1.Cells.Item(198,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(198,4).Value = 97
$ws1.Cells.Item(198,5).Value = "LP1912"
$ws1.Cells.Item(199,1).Value = "12:21:08"
$ws1.Cells.Item(199,2).Value = "12:37"
$ws1.Cells.Item(199,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(199,4).Value = 16
$ws1.Cells.Item(199,5).Value = "LP1912"
$ws1.Cells.Item(200,1).Value = "10:59:49"
$ws1.Cells.Item(200,2).Value = "12:38"
$ws1.Cells.Item(200,3).Value = "17_179 Y 38"
$ws1.Cells.Item(200,4).Value = 99
$ws1.Cells.Item(200,5).Value = "LP1912"
$ws1.Cells.Item(201,1).Value = "11:56:55"
$ws1.Cells.Item(201,2).Value = "12:40"
$ws1.Cells.Item(201,3).Value = "10_OLMOS"
$ws1.Cells.Item(201,4).Value = 44
$ws1.Cells.Item(201,5).Value = "LP1912"
$ws1.Cells.Item(202,1).Value = "11:30:45"
$ws1.Cells.Item(202,2).Value = "12:41"
$ws1.Cells.Item(202,3).Value = "10_OLMOS"
$ws1.Cells.Item(202,4).Value = 71
$ws1.Cells.Item(202,5).Value = "LP1912"
$ws1.Cells.Item(203,1).Value = "12:47:27"
$ws1.Cells.Item(203,2).Value = "12:47"
$ws1.Cells.Item(203,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(203,4).Value = 0
$ws1.Cells.Item(203,5).Value = "LP1912"
$ws1.Cells.Item(204,1).Value = "12:47:27"
$ws1.Cells.Item(204,2).Value = "12:48"
$ws1.Cells.Item(204,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(204,4).Value = 1
$ws1.Cells.Item(204,5).Value = "LP1912"
$ws1.Cells.Item(205,1).Value = "10:59:49"
$ws1.Cells.Item(205,2).Value = "12:48"
$ws1.Cells.Item(205,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(205,4).Value = 109
$ws1.Cells.Item(205,5).Value = "LP1912"
$ws1.Cells.Item(206,1).Value = "12:21:08"
$ws1.Cells.Item(206,2).Value = "12:49"
$ws1.Cells.Item(206,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(206,4).Value = 28
$ws1.Cells.Item(206,5).Value = "LP1912"
$ws1.Cells.Item(207,1).Value = "12:21:08"
$ws1.Cells.Item(207,2).Value = "12:55"
$ws1.Cells.Item(207,3).Value = "10_OLMOS"
$ws1.Cells.Item(207,4).Value = 34
$ws1.Cells.Item(207,5).Value = "LP1912"
$ws1.Cells.Item(208,1).Value = "12:59:47"
$ws1.Cells.Item(208,2).Value = "13:00"
$ws1.Cells.Item(208,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(208,4).Value = 1
$ws1.Cells.Item(208,5).Value = "LP1912"
$ws1.Cells.Item(209,1).Value = "11:30:45"
$ws1.Cells.Item(209,2).Value = "13:01"
$ws1.Cells.Item(209,3).Value = "17_ROMERO"
$ws1.Cells.Item(209,4).Value = 91
$ws1.Cells.Item(209,5).Value = "LP1912"
$ws1.Cells.Item(210,1).Value = "12:47:27"
$ws1.Cells.Item(210,2).Value = "13:02"
$ws1.Cells.Item(210,3).Value = "15_ABASTO"
$ws1.Cells.Item(210,4).Value = 15
$ws1.Cells.Item(210,5).Value = "LP1912"
$ws1.Cells.Item(211,1).Value = "12:21:08"
$ws1.Cells.Item(211,2).Value = "13:03"
$ws1.Cells.Item(211,3).Value = "14_ABASTO"
$ws1.Cells.Item(211,4).Value = 42
$ws1.Cells.Item(211,5).Value = "LP1912"
$ws1.Cells.Item(212,1).Value = "12:47:27"
$ws1.Cells.Item(212,2).Value = "13:04"
$ws1.Cells.Item(212,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(212,4).Value = 17
$ws1.Cells.Item(212,5).Value = "LP1912"
$ws1.Cells.Item(213,1).Value = "12:59:47"
$ws1.Cells.Item(213,2).Value = "13:05"
$ws1.Cells.Item(213,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(213,4).Value = 6
$ws1.Cells.Item(213,5).Value = "LP1912"
$ws1.Cells.Item(214,1).Value = "11:30:45"
$ws1.Cells.Item(214,2).Value = "13:06"
$ws1.Cells.Item(214,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(214,4).Value = 96
$ws1.Cells.Item(214,5).Value = "LP1912"
$ws1.Cells.Item(215,1).Value = "12:21:08"
$ws1.Cells.Item(215,2).Value = "13:07"
$ws1.Cells.Item(215,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(215,4).Value = 46
$ws1.Cells.Item(215,5).Value = "LP1912"
$ws1.Cells.Item(216,1).Value = "11:30:45"
$ws1.Cells.Item(216,2).Value = "13:07"
$ws1.Cells.Item(216,3).Value = "10_OLMOS"
$ws1.Cells.Item(216,4).Value = 97
$ws1.Cells.Item(216,5).Value = "LP1912"
$ws1.Cells.Item(217,1).Value = "12:21:08"
$ws1.Cells.Item(217,2).Value = "13:08"
$ws1.Cells.Item(217,3).Value = "10_OLMOS"
$ws1.Cells.Item(217,4).Value = 47
$ws1.Cells.Item(217,5).Value = "LP1912"
$ws1.Cells.Item(218,1).Value = "11:30:45"
$ws1.Cells.Item(218,2).Value = "13:13"
$ws1.Cells.Item(218,3).Value = "215D_EL PATO"
$ws1.Cells.Item(218,4).Value = 103
$ws1.Cells.Item(218,5).Value = "LP1912"
$ws1.Cells.Item(219,1).Value = "12:21:08"
$ws1.Cells.Item(219,2).Value = "13:14"
$ws1.Cells.Item(219,3).Value = "215D_EL PATO"
$ws1.Cells.Item(219,4).Value = 53
$ws1.Cells.Item(219,5).Value = "LP1912"
$ws1.Cells.Item(220,1).Value = "12:47:27"
$ws1.Cells.Item(220,2).Value = "13:14"
$ws1.Cells.Item(220,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(220,4).Value = 27
$ws1.Cells.Item(220,5).Value = "LP1912"
$ws1.Cells.Item(221,1).Value = "11:56:55"
$ws1.Cells.Item(221,2).Value = "13:20"
$ws1.Cells.Item(221,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(221,4).Value = 84
$ws1.Cells.Item(221,5).Value = "LP1912"
$ws1.Cells.Item(222,1).Value = "11:30:45"
$ws1.Cells.Item(222,2).Value = "13:21"
$ws1.Cells.Item(222,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(222,4).Value = 111
$ws1.Cells.Item(222,5).Value = "LP1912"
$ws1.Cells.Item(223,1).Value = "11:30:45"
$ws1.Cells.Item(223,2).Value = "13:25"
$ws1.Cells.Item(223,3).Value = "10_OLMOS"
$ws1.Cells.Item(223,4).Value = 115
$ws1.Cells.Item(223,5).Value = "LP1912"
$ws1.Cells.Item(224,1).Value = "11:30:45"
$ws1.Cells.Item(224,2).Value = "13:26"
$ws1.Cells.Item(224,3).Value = "15_ABASTO"
$ws1.Cells.Item(224,4).Value = 116
$ws1.Cells.Item(224,5).Value = "LP1912"
$ws1.Cells.Item(225,1).Value = "11:30:45"
$ws1.Cells.Item(225,2).Value = "13:26"
$ws1.Cells.Item(225,3).Value = "14_ABASTO"
$ws1.Cells.Item(225,4).Value = 116
$ws1.Cells.Item(225,5).Value = "LP1912"
$ws1.Cells.Item(226,1).Value = "11:56:55"
$ws1.Cells.Item(226,2).Value = "13:27"
$ws1.Cells.Item(226,3).Value = "10_OLMOS"
$ws1.Cells.Item(226,4).Value = 91
$ws1.Cells.Item(226,5).Value = "LP1912"
$ws1.Cells.Item(227,1).Value = "12:21:08"
$ws1.Cells.Item(227,2).Value = "13:27"
$ws1.Cells.Item(227,3).Value = "14_ABASTO"
$ws1.Cells.Item(227,4).Value = 66
$ws1.Cells.Item(227,5).Value = "LP1912"
$ws1.Cells.Item(228,1).Value = "12:21:08"
$ws1.Cells.Item(228,2).Value = "13:28"
$ws1.Cells.Item(228,3).Value = "10_OLMOS"
$ws1.Cells.Item(228,4).Value = 67
$ws1.Cells.Item(228,5).Value = "LP1912"
$ws1.Cells.Item(229,1).Value = "12:47:27"
$ws1.Cells.Item(229,2).Value = "13:31"
$ws1.Cells.Item(229,3).Value = "10_OLMOS"
$ws1.Cells.Item(229,4).Value = 44
$ws1.Cells.Item(229,5).Value = "LP1912"
$ws1.Cells.Item(230,1).Value = "12:47:27"
$ws1.Cells.Item(230,2).Value = "13:32"
$ws1.Cells.Item(230,3).Value = "10_OLMOS"
$ws1.Cells.Item(230,4).Value = 45
$ws1.Cells.Item(230,5).Value = "LP1912"
$ws1.Cells.Item(231,1).Value = "13:33:42"
$ws1.Cells.Item(231,2).Value = "13:33"
$ws1.Cells.Item(231,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(231,4).Value = 0
$ws1.Cells.Item(231,5).Value = "LP1912"
$ws1.Cells.Item(232,1).Value = "12:59:47"
$ws1.Cells.Item(232,2).Value = "13:33"
$ws1.Cells.Item(232,3).Value = "10_OLMOS"
$ws1.Cells.Item(232,4).Value = 34
$ws1.Cells.Item(232,5).Value = "LP1912"
$ws1.Cells.Item(233,1).Value = "13:33:42"
$ws1.Cells.Item(233,2).Value = "13:34"
$ws1.Cells.Item(233,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(233,4).Value = 1
$ws1.Cells.Item(233,5).Value = "LP1912"
$ws1.Cells.Item(234,1).Value = "13:33:42"
$ws1.Cells.Item(234,2).Value = "13:34"
$ws1.Cells.Item(234,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(234,4).Value = 1
$ws1.Cells.Item(234,5).Value = "LP1912"
$ws1.Cells.Item(235,1).Value = "11:56:55"
$ws1.Cells.Item(235,2).Value = "13:36"
$ws1.Cells.Item(235,3).Value = "15_ABASTO"
$ws1.Cells.Item(235,4).Value = 100
$ws1.Cells.Item(235,5).Value = "LP1912"
$ws1.Cells.Item(236,1).Value = "13:33:42"
$ws1.Cells.Item(236,2).Value = "13:38"
$ws1.Cells.Item(236,3).Value = "14_ABASTO"
$ws1.Cells.Item(236,4).Value = 5
$ws1.Cells.Item(236,5).Value = "LP1912"
$ws1.Cells.Item(237,1).Value = "11:56:55"
$ws1.Cells.Item(237,2).Value = "13:46"
$ws1.Cells.Item(237,3).Value = "17_ROMERO"
$ws1.Cells.Item(237,4).Value = 110
$ws1.Cells.Item(237,5).Value = "LP1912"
$ws1.Cells.Item(238,1).Value = "12:59:47"
$ws1.Cells.Item(238,2).Value = "13:50"
$ws1.Cells.Item(238,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(238,4).Value = 51
$ws1.Cells.Item(238,5).Value = "LP1912"
$ws1.Cells.Item(239,1).Value = "11:56:55"
$ws1.Cells.Item(239,2).Value = "13:50"
$ws1.Cells.Item(239,3).Value = "215A_EL PATO"
$ws1.Cells.Item(239,4).Value = 114
$ws1.Cells.Item(239,5).Value = "LP1912"
$ws1.Cells.Item(240,1).Value = "12:21:08"
$ws1.Cells.Item(240,2).Value = "13:51"
$ws1.Cells.Item(240,3).Value = "215A_EL PATO"
$ws1.Cells.Item(240,4).Value = 90
$ws1.Cells.Item(240,5).Value = "LP1912"
$ws1.Cells.Item(241,1).Value = "11:56:55"
$ws1.Cells.Item(241,2).Value = "13:55"
$ws1.Cells.Item(241,3).Value = "225_GOMEZ"
$ws1.Cells.Item(241,4).Value = 119
$ws1.Cells.Item(241,5).Value = "LP1912"
$ws1.Cells.Item(242,1).Value = "12:21:08"
$ws1.Cells.Item(242,2).Value = "13:56"
$ws1.Cells.Item(242,3).Value = "225_GOMEZ"
$ws1.Cells.Item(242,4).Value = 95
$ws1.Cells.Item(242,5).Value = "LP1912"
$ws1.Cells.Item(243,1).Value = "12:59:47"
$ws1.Cells.Item(243,2).Value = "13:56"
$ws1.Cells.Item(243,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(243,4).Value = 57
$ws1.Cells.Item(243,5).Value = "LP1912"
$ws1.Cells.Item(244,1).Value = "12:47:27"
$ws1.Cells.Item(244,2).Value = "13:58"
$ws1.Cells.Item(244,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(244,4).Value = 71
$ws1.Cells.Item(244,5).Value = "LP1912"
$ws1.Cells.Item(245,1).Value = "13:59:06"
$ws1.Cells.Item(245,2).Value = "13:59"
$ws1.Cells.Item(245,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(245,4).Value = 0
$ws1.Cells.Item(245,5).Value = "LP1912"
$ws1.Cells.Item(246,1).Value = "13:59:06"
$ws1.Cells.Item(246,2).Value = "14:00"
$ws1.Cells.Item(246,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(246,4).Value = 1
$ws1.Cells.Item(246,5).Value = "LP1912"
$ws1.Cells.Item(247,1).Value = "12:21:08"
$ws1.Cells.Item(247,2).Value = "14:00"
$ws1.Cells.Item(247,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(247,4).Value = 99
$ws1.Cells.Item(247,5).Value = "LP1912"
$ws1.Cells.Item(248,1).Value = "12:21:08"
$ws1.Cells.Item(248,2).Value = "14:04"
$ws1.Cells.Item(248,3).Value = "17_ROMERO"
$ws1.Cells.Item(248,4).Value = 103
$ws1.Cells.Item(248,5).Value = "LP1912"
$ws1.Cells.Item(249,1).Value = "13:33:42"
$ws1.Cells.Item(249,2).Value = "14:04"
$ws1.Cells.Item(249,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(249,4).Value = 31
$ws1.Cells.Item(249,5).Value = "LP1912"
$ws1.Cells.Item(250,1).Value = "13:59:06"
$ws1.Cells.Item(250,2).Value = "14:05"
$ws1.Cells.Item(250,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(250,4).Value = 6
$ws1.Cells.Item(250,5).Value = "LP1912"
$ws1.Cells.Item(251,1).Value = "12:21:08"
$ws1.Cells.Item(251,2).Value = "14:08"
$ws1.Cells.Item(251,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(251,4).Value = 107
$ws1.Cells.Item(251,5).Value = "LP1912"
$ws1.Cells.Item(252,1).Value = "12:59:47"
$ws1.Cells.Item(252,2).Value = "14:11"
$ws1.Cells.Item(252,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(252,4).Value = 72
$ws1.Cells.Item(252,5).Value = "LP1912"
$ws1.Cells.Item(253,1).Value = "13:33:42"
$ws1.Cells.Item(253,2).Value = "14:12"
$ws1.Cells.Item(253,3).Value = "15_ABASTO"
$ws1.Cells.Item(253,4).Value = 39
$ws1.Cells.Item(253,5).Value = "LP1912"
$ws1.Cells.Item(254,1).Value = "12:47:27"
$ws1.Cells.Item(254,2).Value = "14:16"
$ws1.Cells.Item(254,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(254,4).Value = 89
$ws1.Cells.Item(254,5).Value = "LP1912"
$ws1.Cells.Item(255,1).Value = "12:21:08"
$ws1.Cells.Item(255,2).Value = "14:17"
$ws1.Cells.Item(255,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(255,4).Value = 116
$ws1.Cells.Item(255,5).Value = "LP1912"
$ws1.Cells.Item(256,1).Value = "12:59:47"
$ws1.Cells.Item(256,2).Value = "14:19"
$ws1.Cells.Item(256,3).Value = "215C_EL PATO"
$ws1.Cells.Item(256,4).Value = 80
$ws1.Cells.Item(256,5).Value = "LP1912"
$ws1.Cells.Item(257,1).Value = "12:21:08"
$ws1.Cells.Item(257,2).Value = "14:20"
$ws1.Cells.Item(257,3).Value = "215C_EL PATO"
$ws1.Cells.Item(257,4).Value = 119
$ws1.Cells.Item(257,5).Value = "LP1912"
$ws1.Cells.Item(258,1).Value = "12:47:27"
$ws1.Cells.Item(258,2).Value = "14:21"
$ws1.Cells.Item(258,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(258,4).Value = 94
$ws1.Cells.Item(258,5).Value = "LP1912"
$ws1.Cells.Item(259,1).Value = "14:24:16"
$ws1.Cells.Item(259,2).Value = "14:25"
$ws1.Cells.Item(259,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(259,4).Value = 1
$ws1.Cells.Item(259,5).Value = "LP1912"
$ws1.Cells.Item(260,1).Value = "13:59:06"
$ws1.Cells.Item(260,2).Value = "14:28"
$ws1.Cells.Item(260,3).Value = "15_ABASTO"
$ws1.Cells.Item(260,4).Value = 29
$ws1.Cells.Item(260,5).Value = "LP1912"
$ws1.Cells.Item(261,1).Value = "14:24:16"
$ws1.Cells.Item(261,2).Value = "14:35"
$ws1.Cells.Item(261,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(261,4).Value = 11
$ws1.Cells.Item(261,5).Value = "LP1912"
$ws1.Cells.Item(262,1).Value = "14:24:16"
$ws1.Cells.Item(262,2).Value = "14:44"
$ws1.Cells.Item(262,3).Value = "15_ABASTO"
$ws1.Cells.Item(262,4).Value = 20
$ws1.Cells.Item(262,5).Value = "LP1912"
$ws1.Cells.Item(263,1).Value = "13:33:42"
$ws1.Cells.Item(263,2).Value = "14:44"
$ws1.Cells.Item(263,3).Value = "14_ABASTO"
$ws1.Cells.Item(263,4).Value = 71
$ws1.Cells.Item(263,5).Value = "LP1912"
$ws1.Cells.Item(264,1).Value = "14:45:17"
$ws1.Cells.Item(264,2).Value = "14:45"
$ws1.Cells.Item(264,3).Value = "15_ABASTO"
$ws1.Cells.Item(264,4).Value = 0
$ws1.Cells.Item(264,5).Value = "LP1912"
$ws1.Cells.Item(265,1).Value = "12:47:27"
$ws1.Cells.Item(265,2).Value = "14:45"
$ws1.Cells.Item(265,3).Value = "14_ABASTO"
$ws1.Cells.Item(265,4).Value = 118
$ws1.Cells.Item(265,5).Value = "LP1912"
$ws1.Cells.Item(266,1).Value = "14:45:17"
$ws1.Cells.Item(266,2).Value = "14:46"
$ws1.Cells.Item(266,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(266,4).Value = 1
$ws1.Cells.Item(266,5).Value = "LP1912"
$ws1.Cells.Item(267,1).Value = "14:56:20"
$ws1.Cells.Item(267,2).Value = "14:56"
$ws1.Cells.Item(267,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(267,4).Value = 0
$ws1.Cells.Item(267,5).Value = "LP1912"
$ws1.Cells.Item(268,1).Value = "12:59:47"
$ws1.Cells.Item(268,2).Value = "14:56"
$ws1.Cells.Item(268,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(268,4).Value = 117
$ws1.Cells.Item(268,5).Value = "LP1912"
$ws1.Cells.Item(269,1).Value = "13:59:06"
$ws1.Cells.Item(269,2).Value = "14:57"
$ws1.Cells.Item(269,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(269,4).Value = 58
$ws1.Cells.Item(269,5).Value = "LP1912"
$ws1.Cells.Item(270,1).Value = "12:59:47"
$ws1.Cells.Item(270,2).Value = "14:58"
$ws1.Cells.Item(270,3).Value = "215B_EL PATO"
$ws1.Cells.Item(270,4).Value = 119
$ws1.Cells.Item(270,5).Value = "LP1912"
$ws1.Cells.Item(271,1).Value = "13:33:42"
$ws1.Cells.Item(271,2).Value = "15:00"
$ws1.Cells.Item(271,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(271,4).Value = 87
$ws1.Cells.Item(271,5).Value = "LP1912"
$ws1.Cells.Item(272,1).Value = "13:33:42"
$ws1.Cells.Item(272,2).Value = "15:05"
$ws1.Cells.Item(272,3).Value = "10_OLMOS"
$ws1.Cells.Item(272,4).Value = 92
$ws1.Cells.Item(272,5).Value = "LP1912"
$ws1.Cells.Item(273,1).Value = "14:45:17"
$ws1.Cells.Item(273,2).Value = "15:05"
$ws1.Cells.Item(273,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(273,4).Value = 20
$ws1.Cells.Item(273,5).Value = "LP1912"
$ws1.Cells.Item(274,1).Value = "13:59:06"
$ws1.Cells.Item(274,2).Value = "15:10"
$ws1.Cells.Item(274,3).Value = "17_ROMERO"
$ws1.Cells.Item(274,4).Value = 71
$ws1.Cells.Item(274,5).Value = "LP1912"
$ws1.Cells.Item(275,1).Value = "13:33:42"
$ws1.Cells.Item(275,2).Value = "15:13"
$ws1.Cells.Item(275,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(275,4).Value = 100
$ws1.Cells.Item(275,5).Value = "LP1912"
$ws1.Cells.Item(276,1).Value = "13:59:06"
$ws1.Cells.Item(276,2).Value = "15:14"
$ws1.Cells.Item(276,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(276,4).Value = 75
$ws1.Cells.Item(276,5).Value = "LP1912"
$ws1.Cells.Item(277,1).Value = "14:56:20"
$ws1.Cells.Item(277,2).Value = "15:17"
$ws1.Cells.Item(277,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(277,4).Value = 21
$ws1.Cells.Item(277,5).Value = "LP1912"
$ws1.Cells.Item(278,1).Value = "13:33:42"
$ws1.Cells.Item(278,2).Value = "15:17"
$ws1.Cells.Item(278,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(278,4).Value = 104
$ws1.Cells.Item(278,5).Value = "LP1912"
$ws1.Cells.Item(279,1).Value = "13:59:06"
$ws1.Cells.Item(279,2).Value = "15:18"
$ws1.Cells.Item(279,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(279,4).Value = 79
$ws1.Cells.Item(279,5).Value = "LP1912"
$ws1.Cells.Item(280,1).Value = "14:56:20"
$ws1.Cells.Item(280,2).Value = "15:20"
$ws1.Cells.Item(280,3).Value = "15_ABASTO"
$ws1.Cells.Item(280,4).Value = 24
$ws1.Cells.Item(280,5).Value = "LP1912"
$ws1.Cells.Item(281,1).Value = "14:24:16"
$ws1.Cells.Item(281,2).Value = "15:21"
$ws1.Cells.Item(281,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(281,4).Value = 57
$ws1.Cells.Item(281,5).Value = "LP1912"
$ws1.Cells.Item(282,1).Value = "15:22:17"
$ws1.Cells.Item(282,2).Value = "15:22"
$ws1.Cells.Item(282,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(282,4).Value = 0
$ws1.Cells.Item(282,5).Value = "LP1912"
$ws1.Cells.Item(283,1).Value = "15:22:17"
$ws1.Cells.Item(283,2).Value = "15:22"
$ws1.Cells.Item(283,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(283,4).Value = 0
$ws1.Cells.Item(283,5).Value = "LP1912"
$ws1.Cells.Item(284,1).Value = "14:24:16"
$ws1.Cells.Item(284,2).Value = "15:32"
$ws1.Cells.Item(284,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(284,4).Value = 68
$ws1.Cells.Item(284,5).Value = "LP1912"
$ws1.Cells.Item(285,1).Value = "13:59:06"
$ws1.Cells.Item(285,2).Value = "15:35"
$ws1.Cells.Item(285,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(285,4).Value = 96
$ws1.Cells.Item(285,5).Value = "LP1912"
$ws1.Cells.Item(286,1).Value = "13:59:06"
$ws1.Cells.Item(286,2).Value = "15:37"
$ws1.Cells.Item(286,3).Value = "10_OLMOS"
$ws1.Cells.Item(286,4).Value = 98
$ws1.Cells.Item(286,5).Value = "LP1912"
$ws1.Cells.Item(287,1).Value = "14:56:20"
$ws1.Cells.Item(287,2).Value = "15:38"
$ws1.Cells.Item(287,3).Value = "10_OLMOS"
$ws1.Cells.Item(287,4).Value = 42
$ws1.Cells.Item(287,5).Value = "LP1912"
$ws1.Cells.Item(288,1).Value = "14:45:17"
$ws1.Cells.Item(288,2).Value = "15:38"
$ws1.Cells.Item(288,3).Value = "215A_EL PATO"
$ws1.Cells.Item(288,4).Value = 53
$ws1.Cells.Item(288,5).Value = "LP1912"
$ws1.Cells.Item(289,1).Value = "14:24:16"
$ws1.Cells.Item(289,2).Value = "15:38"
$ws1.Cells.Item(289,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(289,4).Value = 74
$ws1.Cells.Item(289,5).Value = "LP1912"
$ws1.Cells.Item(290,1).Value = "13:59:06"
$ws1.Cells.Item(290,2).Value = "15:39"
$ws1.Cells.Item(290,3).Value = "215A_EL PATO"
$ws1.Cells.Item(290,4).Value = 100
$ws1.Cells.Item(290,5).Value = "LP1912"
$ws1.Cells.Item(291,1).Value = "14:56:20"
$ws1.Cells.Item(291,2).Value = "15:45"
$ws1.Cells.Item(291,3).Value = "14_ABASTO"
$ws1.Cells.Item(291,4).Value = 49
$ws1.Cells.Item(291,5).Value = "LP1912"
$ws1.Cells.Item(292,1).Value = "14:24:16"
$ws1.Cells.Item(292,2).Value = "15:46"
$ws1.Cells.Item(292,3).Value = "14_ABASTO"
$ws1.Cells.Item(292,4).Value = 82
$ws1.Cells.Item(292,5).Value = "LP1912"
$ws1.Cells.Item(293,1).Value = "14:56:20"
$ws1.Cells.Item(293,2).Value = "15:46"
$ws1.Cells.Item(293,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(293,4).Value = 50
$ws1.Cells.Item(293,5).Value = "LP1912"
$ws1.Cells.Item(294,1).Value = "13:59:06"
$ws1.Cells.Item(294,2).Value = "15:47"
$ws1.Cells.Item(294,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(294,4).Value = 108
$ws1.Cells.Item(294,5).Value = "LP1912"
$ws1.Cells.Item(295,1).Value = "13:59:06"
$ws1.Cells.Item(295,2).Value = "15:48"
$ws1.Cells.Item(295,3).Value = "14_ABASTO"
$ws1.Cells.Item(295,4).Value = 109
$ws1.Cells.Item(295,5).Value = "LP1912"
$ws1.Cells.Item(296,1).Value = "14:56:20"
$ws1.Cells.Item(296,2).Value = "15:53"
$ws1.Cells.Item(296,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(296,4).Value = 57
$ws1.Cells.Item(296,5).Value = "LP1912"
$ws1.Cells.Item(297,1).Value = "13:59:06"
$ws1.Cells.Item(297,2).Value = "15:54"
$ws1.Cells.Item(297,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(297,4).Value = 115
$ws1.Cells.Item(297,5).Value = "LP1912"
$ws1.Cells.Item(298,1).Value = "15:53:28"
$ws1.Cells.Item(298,2).Value = "15:54"
$ws1.Cells.Item(298,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(298,4).Value = 1
$ws1.Cells.Item(298,5).Value = "LP1912"
$ws1.Cells.Item(299,1).Value = "15:22:17"
$ws1.Cells.Item(299,2).Value = "15:55"
$ws1.Cells.Item(299,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(299,4).Value = 33
$ws1.Cells.Item(299,5).Value = "LP1912"
$ws1.Cells.Item(300,1).Value = "14:24:16"
$ws1.Cells.Item(300,2).Value = "15:56"
$ws1.Cells.Item(300,3).Value = "17_ROMERO"
$ws1.Cells.Item(300,4).Value = 92
$ws1.Cells.Item(300,5).Value = "LP1912"
$ws1.Cells.Item(301,1).Value = "15:53:28"
$ws1.Cells.Item(301,2).Value = "15:56"
$ws1.Cells.Item(301,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(301,4).Value = 3
$ws1.Cells.Item(301,5).Value = "LP1912"
$ws1.Cells.Item(302,1).Value = "13:59:06"
$ws1.Cells.Item(302,2).Value = "15:57"
$ws1.Cells.Item(302,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(302,4).Value = 118
$ws1.Cells.Item(302,5).Value = "LP1912"
$ws1.Cells.Item(303,1).Value = "15:22:17"
$ws1.Cells.Item(303,2).Value = "16:01"
$ws1.Cells.Item(303,3).Value = "10_OLMOS"
$ws1.Cells.Item(303,4).Value = 39
$ws1.Cells.Item(303,5).Value = "LP1912"
$ws1.Cells.Item(304,1).Value = "15:53:28"
$ws1.Cells.Item(304,2).Value = "16:02"
$ws1.Cells.Item(304,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(304,4).Value = 9
$ws1.Cells.Item(304,5).Value = "LP1912"
$ws1.Cells.Item(305,1).Value = "15:53:28"
$ws1.Cells.Item(305,2).Value = "16:04"
$ws1.Cells.Item(305,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(305,4).Value = 11
$ws1.Cells.Item(305,5).Value = "LP1912"
$ws1.Cells.Item(306,1).Value = "15:22:17"
$ws1.Cells.Item(306,2).Value = "16:05"
$ws1.Cells.Item(306,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(306,4).Value = 43
$ws1.Cells.Item(306,5).Value = "LP1912"
$ws1.Cells.Item(307,1).Value = "14:56:20"
$ws1.Cells.Item(307,2).Value = "16:08"
$ws1.Cells.Item(307,3).Value = "14_ABASTO"
$ws1.Cells.Item(307,4).Value = 72
$ws1.Cells.Item(307,5).Value = "LP1912"
$ws1.Cells.Item(308,1).Value = "14:45:17"
$ws1.Cells.Item(308,2).Value = "16:09"
$ws1.Cells.Item(308,3).Value = "14_ABASTO"
$ws1.Cells.Item(308,4).Value = 84
$ws1.Cells.Item(308,5).Value = "LP1912"
$ws1.Cells.Item(309,1).Value = "14:24:16"
$ws1.Cells.Item(309,2).Value = "16:15"
$ws1.Cells.Item(309,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(309,4).Value = 111
$ws1.Cells.Item(309,5).Value = "LP1912"
$ws1.Cells.Item(310,1).Value = "16:13:37"
$ws1.Cells.Item(310,2).Value = "16:15"
$ws1.Cells.Item(310,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(310,4).Value = 2
$ws1.Cells.Item(310,5).Value = "LP1912"
$ws1.Cells.Item(311,1).Value = "16:13:37"
$ws1.Cells.Item(311,2).Value = "16:19"
$ws1.Cells.Item(311,3).Value = "215C_EL PATO"
$ws1.Cells.Item(311,4).Value = 6
$ws1.Cells.Item(311,5).Value = "LP1912"
$ws1.Cells.Item(312,1).Value = "14:24:16"
$ws1.Cells.Item(312,2).Value = "16:20"
$ws1.Cells.Item(312,3).Value = "215C_EL PATO"
$ws1.Cells.Item(312,4).Value = 116
$ws1.Cells.Item(312,5).Value = "LP1912"
$ws1.Cells.Item(313,1).Value = "14:24:16"
$ws1.Cells.Item(313,2).Value = "16:21"
$ws1.Cells.Item(313,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(313,4).Value = 117
$ws1.Cells.Item(313,5).Value = "LP1912"
$ws1.Cells.Item(314,1).Value = "16:13:37"
$ws1.Cells.Item(314,2).Value = "16:26"
$ws1.Cells.Item(314,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(314,4).Value = 13
$ws1.Cells.Item(314,5).Value = "LP1912"
$ws1.Cells.Item(315,1).Value = "15:53:28"
$ws1.Cells.Item(315,2).Value = "16:29"
$ws1.Cells.Item(315,3).Value = "10_OLMOS"
$ws1.Cells.Item(315,4).Value = 36
$ws1.Cells.Item(315,5).Value = "LP1912"
$ws1.Cells.Item(316,1).Value = "14:45:17"
$ws1.Cells.Item(316,2).Value = "16:30"
$ws1.Cells.Item(316,3).Value = "15_ABASTO"
$ws1.Cells.Item(316,4).Value = 105
$ws1.Cells.Item(316,5).Value = "LP1912"
$ws1.Cells.Item(317,1).Value = "16:31:51"
$ws1.Cells.Item(317,2).Value = "16:31"
$ws1.Cells.Item(317,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(317,4).Value = 0
$ws1.Cells.Item(317,5).Value = "LP1912"
$ws1.Cells.Item(318,1).Value = "15:22:17"
$ws1.Cells.Item(318,2).Value = "16:32"
$ws1.Cells.Item(318,3).Value = "14_ABASTO"
$ws1.Cells.Item(318,4).Value = 70
$ws1.Cells.Item(318,5).Value = "LP1912"
$ws1.Cells.Item(319,1).Value = "15:53:28"
$ws1.Cells.Item(319,2).Value = "16:34"
$ws1.Cells.Item(319,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(319,4).Value = 41
$ws1.Cells.Item(319,5).Value = "LP1912"
$ws1.Cells.Item(320,1).Value = "15:53:28"
$ws1.Cells.Item(320,2).Value = "16:36"
$ws1.Cells.Item(320,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(320,4).Value = 43
$ws1.Cells.Item(320,5).Value = "LP1912"
$ws1.Cells.Item(321,1).Value = "15:22:17"
$ws1.Cells.Item(321,2).Value = "16:37"
$ws1.Cells.Item(321,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(321,4).Value = 75
$ws1.Cells.Item(321,5).Value = "LP1912"
$ws1.Cells.Item(322,1).Value = "15:22:17"
$ws1.Cells.Item(322,2).Value = "16:40"
$ws1.Cells.Item(322,3).Value = "17_ROMERO"
$ws1.Cells.Item(322,4).Value = 78
$ws1.Cells.Item(322,5).Value = "LP1912"
$ws1.Cells.Item(323,1).Value = "14:56:20"
$ws1.Cells.Item(323,2).Value = "16:42"
$ws1.Cells.Item(323,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(323,4).Value = 106
$ws1.Cells.Item(323,5).Value = "LP1912"
$ws1.Cells.Item(324,1).Value = "16:31:51"
$ws1.Cells.Item(324,2).Value = "16:42"
$ws1.Cells.Item(324,3).Value = "225_GOMEZ"
$ws1.Cells.Item(324,4).Value = 11
$ws1.Cells.Item(324,5).Value = "LP1912"
$ws1.Cells.Item(325,1).Value = "14:45:17"
$ws1.Cells.Item(325,2).Value = "16:43"
$ws1.Cells.Item(325,3).Value = "225_GOMEZ"
$ws1.Cells.Item(325,4).Value = 118
$ws1.Cells.Item(325,5).Value = "LP1912"
$ws1.Cells.Item(326,1).Value = "14:45:17"
$ws1.Cells.Item(326,2).Value = "16:43"
$ws1.Cells.Item(326,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(326,4).Value = 118
$ws1.Cells.Item(326,5).Value = "LP1912"
$ws1.Cells.Item(327,1).Value = "16:45:31"
$ws1.Cells.Item(327,2).Value = "16:45"
$ws1.Cells.Item(327,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(327,4).Value = 0
$ws1.Cells.Item(327,5).Value = "LP1912"
$ws1.Cells.Item(328,1).Value = "15:22:17"
$ws1.Cells.Item(328,2).Value = "16:48"
$ws1.Cells.Item(328,3).Value = "15_ABASTO"
$ws1.Cells.Item(328,4).Value = 86
$ws1.Cells.Item(328,5).Value = "LP1912"
$ws1.Cells.Item(329,1).Value = "15:53:28"
$ws1.Cells.Item(329,2).Value = "16:50"
$ws1.Cells.Item(329,3).Value = "14_ABASTO"
$ws1.Cells.Item(329,4).Value = 57
$ws1.Cells.Item(329,5).Value = "LP1912"
$ws1.Cells.Item(330,1).Value = "16:52:36"
$ws1.Cells.Item(330,2).Value = "16:52"
$ws1.Cells.Item(330,3).Value = "10_OLMOS"
$ws1.Cells.Item(330,4).Value = 0
$ws1.Cells.Item(330,5).Value = "LP1912"
$ws1.Cells.Item(331,1).Value = "16:52:36"
$ws1.Cells.Item(331,2).Value = "16:53"
$ws1.Cells.Item(331,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(331,4).Value = 1
$ws1.Cells.Item(331,5).Value = "LP1912"
$ws1.Cells.Item(332,1).Value = "16:31:51"
$ws1.Cells.Item(332,2).Value = "16:56"
$ws1.Cells.Item(332,3).Value = "10_OLMOS"
$ws1.Cells.Item(332,4).Value = 25
$ws1.Cells.Item(332,5).Value = "LP1912"
$ws1.Cells.Item(333,1).Value = "15:22:17"
$ws1.Cells.Item(333,2).Value = "16:56"
$ws1.Cells.Item(333,3).Value = "17_179 Y 38"
$ws1.Cells.Item(333,4).Value = 94
$ws1.Cells.Item(333,5).Value = "LP1912"
$ws1.Cells.Item(334,1).Value = "16:13:37"
$ws1.Cells.Item(334,2).Value = "16:57"
$ws1.Cells.Item(334,3).Value = "10_OLMOS"
$ws1.Cells.Item(334,4).Value = 44
$ws1.Cells.Item(334,5).Value = "LP1912"
$ws1.Cells.Item(335,1).Value = "16:13:37"
$ws1.Cells.Item(335,2).Value = "17:04"
$ws1.Cells.Item(335,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(335,4).Value = 51
$ws1.Cells.Item(335,5).Value = "LP1912"
$ws1.Cells.Item(336,1).Value = "15:22:17"
$ws1.Cells.Item(336,2).Value = "17:04"
$ws1.Cells.Item(336,3).Value = "215A_EL PATO"
$ws1.Cells.Item(336,4).Value = 102
$ws1.Cells.Item(336,5).Value = "LP1912"
$ws1.Cells.Item(337,1).Value = "16:13:37"
$ws1.Cells.Item(337,2).Value = "17:04"
$ws1.Cells.Item(337,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(337,4).Value = 51
$ws1.Cells.Item(337,5).Value = "LP1912"
$ws1.Cells.Item(338,1).Value = "16:52:36"
$ws1.Cells.Item(338,2).Value = "17:05"
$ws1.Cells.Item(338,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(338,4).Value = 13
$ws1.Cells.Item(338,5).Value = "LP1912"
$ws1.Cells.Item(339,1).Value = "16:45:31"
$ws1.Cells.Item(339,2).Value = "17:06"
$ws1.Cells.Item(339,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(339,4).Value = 21
$ws1.Cells.Item(339,5).Value = "LP1912"
$ws1.Cells.Item(340,1).Value = "16:31:51"
$ws1.Cells.Item(340,2).Value = "17:09"
$ws1.Cells.Item(340,3).Value = "10_OLMOS"
$ws1.Cells.Item(340,4).Value = 38
$ws1.Cells.Item(340,5).Value = "LP1912"
$ws1.Cells.Item(341,1).Value = "16:45:31"
$ws1.Cells.Item(341,2).Value = "17:10"
$ws1.Cells.Item(341,3).Value = "10_OLMOS"
$ws1.Cells.Item(341,4).Value = 25
$ws1.Cells.Item(341,5).Value = "LP1912"
$ws1.Cells.Item(342,1).Value = "17:14:34"
$ws1.Cells.Item(342,2).Value = "17:15"
$ws1.Cells.Item(342,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(342,4).Value = 1
$ws1.Cells.Item(342,5).Value = "LP1912"
$ws1.Cells.Item(343,1).Value = "16:45:31"
$ws1.Cells.Item(343,2).Value = "17:16"
$ws1.Cells.Item(343,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(343,4).Value = 31
$ws1.Cells.Item(343,5).Value = "LP1912"
$ws1.Cells.Item(344,1).Value = "17:14:34"
$ws1.Cells.Item(344,2).Value = "17:18"
$ws1.Cells.Item(344,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(344,4).Value = 4
$ws1.Cells.Item(344,5).Value = "LP1912"
$ws1.Cells.Item(345,1).Value = "16:31:51"
$ws1.Cells.Item(345,2).Value = "17:20"
$ws1.Cells.Item(345,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(345,4).Value = 49
$ws1.Cells.Item(345,5).Value = "LP1912"
$ws1.Cells.Item(346,1).Value = "16:31:51"
$ws1.Cells.Item(346,2).Value = "17:20"
$ws1.Cells.Item(346,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(346,4).Value = 49
$ws1.Cells.Item(346,5).Value = "LP1912"
$ws1.Cells.Item(347,1).Value = "15:53:28"
$ws1.Cells.Item(347,2).Value = "17:21"
$ws1.Cells.Item(347,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(347,4).Value = 88
$ws1.Cells.Item(347,5).Value = "LP1912"
$ws1.Cells.Item(348,1).Value = "17:14:34"
$ws1.Cells.Item(348,2).Value = "17:22"
$ws1.Cells.Item(348,3).Value = "10_OLMOS"
$ws1.Cells.Item(348,4).Value = 8
$ws1.Cells.Item(348,5).Value = "LP1912"
$ws1.Cells.Item(349,1).Value = "15:53:28"
$ws1.Cells.Item(349,2).Value = "17:24"
$ws1.Cells.Item(349,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(349,4).Value = 91
$ws1.Cells.Item(349,5).Value = "LP1912"
$ws1.Cells.Item(350,1).Value = "15:53:28"
$ws1.Cells.Item(350,2).Value = "17:28"
$ws1.Cells.Item(350,3).Value = "14_ABASTO"
$ws1.Cells.Item(350,4).Value = 95
$ws1.Cells.Item(350,5).Value = "LP1912"
$ws1.Cells.Item(351,1).Value = "17:14:34"
$ws1.Cells.Item(351,2).Value = "17:34"
$ws1.Cells.Item(351,3).Value = "10_OLMOS"
$ws1.Cells.Item(351,4).Value = 20
$ws1.Cells.Item(351,5).Value = "LP1912"
$ws1.Cells.Item(352,1).Value = "16:45:31"
$ws1.Cells.Item(352,2).Value = "17:34"
$ws1.Cells.Item(352,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(352,4).Value = 49
$ws1.Cells.Item(352,5).Value = "LP1912"
$ws1.Cells.Item(353,1).Value = "15:53:28"
$ws1.Cells.Item(353,2).Value = "17:36"
$ws1.Cells.Item(353,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(353,4).Value = 103
$ws1.Cells.Item(353,5).Value = "LP1912"
$ws1.Cells.Item(354,1).Value = "17:38:19"
$ws1.Cells.Item(354,2).Value = "17:38"
$ws1.Cells.Item(354,3).Value = "215B_EL PATO"
$ws1.Cells.Item(354,4).Value = 0
$ws1.Cells.Item(354,5).Value = "LP1912"
$ws1.Cells.Item(355,1).Value = "15:53:28"
$ws1.Cells.Item(355,2).Value = "17:38"
$ws1.Cells.Item(355,3).Value = "17_ROMERO"
$ws1.Cells.Item(355,4).Value = 105
$ws1.Cells.Item(355,5).Value = "LP1912"
$ws1.Cells.Item(356,1).Value = "17:38:19"
$ws1.Cells.Item(356,2).Value = "17:39"
$ws1.Cells.Item(356,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(356,4).Value = 1
$ws1.Cells.Item(356,5).Value = "LP1912"
$ws1.Cells.Item(357,1).Value = "17:38:19"
$ws1.Cells.Item(357,2).Value = "17:39"
$ws1.Cells.Item(357,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(357,4).Value = 1
$ws1.Cells.Item(357,5).Value = "LP1912"
$ws1.Cells.Item(358,1).Value = "16:13:37"
$ws1.Cells.Item(358,2).Value = "17:40"
$ws1.Cells.Item(358,3).Value = "17_ROMERO"
$ws1.Cells.Item(358,4).Value = 87
$ws1.Cells.Item(358,5).Value = "LP1912"
$ws1.Cells.Item(359,1).Value = "15:53:28"
$ws1.Cells.Item(359,2).Value = "17:40"
$ws1.Cells.Item(359,3).Value = "215B_EL PATO"
$ws1.Cells.Item(359,4).Value = 107
$ws1.Cells.Item(359,5).Value = "LP1912"
$ws1.Cells.Item(360,1).Value = "16:31:51"
$ws1.Cells.Item(360,2).Value = "17:45"
$ws1.Cells.Item(360,3).Value = "15_ABASTO"
$ws1.Cells.Item(360,4).Value = 74
$ws1.Cells.Item(360,5).Value = "LP1912"
$ws1.Cells.Item(361,1).Value = "17:38:19"
$ws1.Cells.Item(361,2).Value = "17:46"
$ws1.Cells.Item(361,3).Value = "10_OLMOS"
$ws1.Cells.Item(361,4).Value = 8
$ws1.Cells.Item(361,5).Value = "LP1912"
$ws1.Cells.Item(362,1).Value = "15:53:28"
$ws1.Cells.Item(362,2).Value = "17:50"
$ws1.Cells.Item(362,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(362,4).Value = 117
$ws1.Cells.Item(362,5).Value = "LP1912"
$ws1.Cells.Item(363,1).Value = "15:53:28"
$ws1.Cells.Item(363,2).Value = "17:52"
$ws1.Cells.Item(363,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(363,4).Value = 119
$ws1.Cells.Item(363,5).Value = "LP1912"
$ws1.Cells.Item(364,1).Value = "17:38:19"
$ws1.Cells.Item(364,2).Value = "17:52"
$ws1.Cells.Item(364,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(364,4).Value = 14
$ws1.Cells.Item(364,5).Value = "LP1912"
$ws1.Cells.Item(365,1).Value = "17:14:34"
$ws1.Cells.Item(365,2).Value = "18:04"
$ws1.Cells.Item(365,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(365,4).Value = 50
$ws1.Cells.Item(365,5).Value = "LP1912"
$ws1.Cells.Item(366,1).Value = "16:13:37"
$ws1.Cells.Item(366,2).Value = "18:04"
$ws1.Cells.Item(366,3).Value = "17_ROMERO"
$ws1.Cells.Item(366,4).Value = 111
$ws1.Cells.Item(366,5).Value = "LP1912"
$ws1.Cells.Item(367,1).Value = "17:38:19"
$ws1.Cells.Item(367,2).Value = "18:05"
$ws1.Cells.Item(367,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(367,4).Value = 27
$ws1.Cells.Item(367,5).Value = "LP1912"
$ws1.Cells.Item(368,1).Value = "16:52:36"
$ws1.Cells.Item(368,2).Value = "18:08"
$ws1.Cells.Item(368,3).Value = "14_ABASTO"
$ws1.Cells.Item(368,4).Value = 76
$ws1.Cells.Item(368,5).Value = "LP1912"
$ws1.Cells.Item(369,1).Value = "17:38:19"
$ws1.Cells.Item(369,2).Value = "18:16"
$ws1.Cells.Item(369,3).Value = "10_OLMOS"
$ws1.Cells.Item(369,4).Value = 38
$ws1.Cells.Item(369,5).Value = "LP1912"
$ws1.Cells.Item(370,1).Value = "16:31:51"
$ws1.Cells.Item(370,2).Value = "18:20"
$ws1.Cells.Item(370,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(370,4).Value = 109
$ws1.Cells.Item(370,5).Value = "LP1912"
$ws1.Cells.Item(371,1).Value = "16:45:31"
$ws1.Cells.Item(371,2).Value = "18:21"
$ws1.Cells.Item(371,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(371,4).Value = 96
$ws1.Cells.Item(371,5).Value = "LP1912"
$ws1.Cells.Item(372,1).Value = "17:38:19"
$ws1.Cells.Item(372,2).Value = "18:24"
$ws1.Cells.Item(372,3).Value = "14_ABASTO"
$ws1.Cells.Item(372,4).Value = 46
$ws1.Cells.Item(372,5).Value = "LP1912"
$ws1.Cells.Item(373,1).Value = "16:31:51"
$ws1.Cells.Item(373,2).Value = "18:27"
$ws1.Cells.Item(373,3).Value = "215C_EL PATO"
$ws1.Cells.Item(373,4).Value = 116
$ws1.Cells.Item(373,5).Value = "LP1912"
$ws1.Cells.Item(374,1).Value = "16:45:31"
$ws1.Cells.Item(374,2).Value = "18:28"
$ws1.Cells.Item(374,3).Value = "215C_EL PATO"
$ws1.Cells.Item(374,4).Value = 103
$ws1.Cells.Item(374,5).Value = "LP1912"
$ws1.Cells.Item(375,1).Value = "16:45:31"
$ws1.Cells.Item(375,2).Value = "18:32"
$ws1.Cells.Item(375,3).Value = "11X44_ETCHEVERRY"
$ws1.Cells.Item(375,4).Value = 107
$ws1.Cells.Item(375,5).Value = "LP1912"
$ws1.Cells.Item(376,1).Value = "16:45:31"
$ws1.Cells.Item(376,2).Value = "18:40"
$ws1.Cells.Item(376,3).Value = "15_ABASTO"
$ws1.Cells.Item(376,4).Value = 115
$ws1.Cells.Item(376,5).Value = "LP1912"
$ws1.Cells.Item(377,1).Value = "16:52:36"
$ws1.Cells.Item(377,2).Value = "18:48"
$ws1.Cells.Item(377,3).Value = "14X44_ABASTO"
$ws1.Cells.Item(377,4).Value = 116
$ws1.Cells.Item(377,5).Value = "LP1912"
$ws1.Cells.Item(378,1).Value = "17:38:19"
$ws1.Cells.Item(378,2).Value = "18:52"
$ws1.Cells.Item(378,3).Value = "15_ABASTO"
$ws1.Cells.Item(378,4).Value = 74
$ws1.Cells.Item(378,5).Value = "LP1912"
$ws1.Cells.Item(379,1).Value = "17:14:34"
$ws1.Cells.Item(379,2).Value = "18:58"
$ws1.Cells.Item(379,3).Value = "215A_EL PATO"
$ws1.Cells.Item(379,4).Value = 104
$ws1.Cells.Item(379,5).Value = "LP1912"
$ws1.Cells.Item(380,1).Value = "17:14:34"
$ws1.Cells.Item(380,2).Value = "19:04"
$ws1.Cells.Item(380,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(380,4).Value = 110
$ws1.Cells.Item(380,5).Value = "LP1912"
$ws1.Cells.Item(381,1).Value = "17:14:34"
$ws1.Cells.Item(381,2).Value = "19:10"
$ws1.Cells.Item(381,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(381,4).Value = 116
$ws1.Cells.Item(381,5).Value = "LP1912"
$ws1.Cells.Item(382,1).Value = "17:38:19"
$ws1.Cells.Item(382,2).Value = "19:17"
$ws1.Cells.Item(382,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(382,4).Value = 99
$ws1.Cells.Item(382,5).Value = "LP1912"
$ws1.Cells.Item(383,1).Value = "17:38:19"
$ws1.Cells.Item(383,2).Value = "19:21"
$ws1.Cells.Item(383,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(383,4).Value = 103
$ws1.Cells.Item(383,5).Value = "LP1912"
$ws1.Cells.Item(384,1).Value = "17:38:19"
$ws1.Cells.Item(384,2).Value = "19:30"
$ws1.Cells.Item(384,3).Value = "225_GOMEZ"
$ws1.Cells.Item(384,4).Value = 112
$ws1.Cells.Item(384,5).Value = "LP1912"
$ws1.Range("A2").Value = "Última actualización: 17:38:19"
$ws1.Range("A3").Value = "Total filas: 379"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215  (39 -> 40 data rows)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(6,1).Value = "03:45:25"
$ws2.Cells.Item(6,2).Value = "04:46"
$ws2.Cells.Item(6,3).Value = "215A_EL PATO"
$ws2.Cells.Item(6,4).Value = 61
$ws2.Cells.Item(6,5).Value = "LP1912"
$ws2.Cells.Item(7,1).Value = "03:45:25"
$ws2.Cells.Item(7,2).Value = "05:34"
$ws2.Cells.Item(7,3).Value = "215B_EL PATO"
$ws2.Cells.Item(7,4).Value = 109
$ws2.Cells.Item(7,5).Value = "LP1912"
$ws2.Cells.Item(8,1).Value = "04:18:02"
$ws2.Cells.Item(8,2).Value = "05:35"
$ws2.Cells.Item(8,3).Value = "215B_EL PATO"
$ws2.Cells.Item(8,4).Value = 77
$ws2.Cells.Item(8,5).Value = "LP1912"
$ws2.Cells.Item(9,1).Value = "04:56:49"
$ws2.Cells.Item(9,2).Value = "06:11"
$ws2.Cells.Item(9,3).Value = "215A_EL PATO"
$ws2.Cells.Item(9,4).Value = 75
$ws2.Cells.Item(9,5).Value = "LP1912"
$ws2.Cells.Item(10,1).Value = "04:18:02"
$ws2.Cells.Item(10,2).Value = "06:12"
$ws2.Cells.Item(10,3).Value = "215A_EL PATO"
$ws2.Cells.Item(10,4).Value = 114
$ws2.Cells.Item(10,5).Value = "LP1912"
$ws2.Cells.Item(11,1).Value = "04:56:49"
$ws2.Cells.Item(11,2).Value = "06:46"
$ws2.Cells.Item(11,3).Value = "215C_EL PATO"
$ws2.Cells.Item(11,4).Value = 110
$ws2.Cells.Item(11,5).Value = "LP1912"
$ws2.Cells.Item(12,1).Value = "05:26:08"
$ws2.Cells.Item(12,2).Value = "06:47"
$ws2.Cells.Item(12,3).Value = "215C_EL PATO"
$ws2.Cells.Item(12,4).Value = 81
$ws2.Cells.Item(12,5).Value = "LP1912"
$ws2.Cells.Item(13,1).Value = "05:26:08"
$ws2.Cells.Item(13,2).Value = "07:11"
$ws2.Cells.Item(13,3).Value = "215A_EL PATO"
$ws2.Cells.Item(13,4).Value = 105
$ws2.Cells.Item(13,5).Value = "LP1912"
$ws2.Cells.Item(14,1).Value = "06:55:02"
$ws2.Cells.Item(14,2).Value = "07:12"
$ws2.Cells.Item(14,3).Value = "215A_EL PATO"
$ws2.Cells.Item(14,4).Value = 17
$ws2.Cells.Item(14,5).Value = "LP1912"
$ws2.Cells.Item(15,1).Value = "05:55:25"
$ws2.Cells.Item(15,2).Value = "07:51"
$ws2.Cells.Item(15,3).Value = "215D_EL PATO"
$ws2.Cells.Item(15,4).Value = 116
$ws2.Cells.Item(15,5).Value = "LP1912"
$ws2.Cells.Item(16,1).Value = "06:55:02"
$ws2.Cells.Item(16,2).Value = "07:52"
$ws2.Cells.Item(16,3).Value = "215D_EL PATO"
$ws2.Cells.Item(16,4).Value = 57
$ws2.Cells.Item(16,5).Value = "LP1912"
$ws2.Cells.Item(17,1).Value = "06:25:43"
$ws2.Cells.Item(17,2).Value = "08:23"
$ws2.Cells.Item(17,3).Value = "215B_EL PATO"
$ws2.Cells.Item(17,4).Value = 118
$ws2.Cells.Item(17,5).Value = "LP1912"
$ws2.Cells.Item(18,1).Value = "07:19:29"
$ws2.Cells.Item(18,2).Value = "09:01"
$ws2.Cells.Item(18,3).Value = "215A_EL PATO"
$ws2.Cells.Item(18,4).Value = 102
$ws2.Cells.Item(18,5).Value = "LP1912"
$ws2.Cells.Item(19,1).Value = "07:50:16"
$ws2.Cells.Item(19,2).Value = "09:02"
$ws2.Cells.Item(19,3).Value = "215A_EL PATO"
$ws2.Cells.Item(19,4).Value = 72
$ws2.Cells.Item(19,5).Value = "LP1912"
$ws2.Cells.Item(20,1).Value = "07:50:16"
$ws2.Cells.Item(20,2).Value = "09:42"
$ws2.Cells.Item(20,3).Value = "215C_EL PATO"
$ws2.Cells.Item(20,4).Value = 112
$ws2.Cells.Item(20,5).Value = "LP1912"
$ws2.Cells.Item(21,1).Value = "08:56:29"
$ws2.Cells.Item(21,2).Value = "10:26"
$ws2.Cells.Item(21,3).Value = "215A_EL PATO"
$ws2.Cells.Item(21,4).Value = 90
$ws2.Cells.Item(21,5).Value = "LP1912"
$ws2.Cells.Item(22,1).Value = "08:32:09"
$ws2.Cells.Item(22,2).Value = "10:27"
$ws2.Cells.Item(22,3).Value = "215A_EL PATO"
$ws2.Cells.Item(22,4).Value = 115
$ws2.Cells.Item(22,5).Value = "LP1912"
$ws2.Cells.Item(23,1).Value = "10:29:57"
$ws2.Cells.Item(23,2).Value = "11:01"
$ws2.Cells.Item(23,3).Value = "215C_EL PATO"
$ws2.Cells.Item(23,4).Value = 32
$ws2.Cells.Item(23,5).Value = "LP1912"
$ws2.Cells.Item(24,1).Value = "09:35:26"
$ws2.Cells.Item(24,2).Value = "11:02"
$ws2.Cells.Item(24,3).Value = "215C_EL PATO"
$ws2.Cells.Item(24,4).Value = 87
$ws2.Cells.Item(24,5).Value = "LP1912"
$ws2.Cells.Item(25,1).Value = "10:29:57"
$ws2.Cells.Item(25,2).Value = "11:50"
$ws2.Cells.Item(25,3).Value = "215B_EL PATO"
$ws2.Cells.Item(25,4).Value = 81
$ws2.Cells.Item(25,5).Value = "LP1912"
$ws2.Cells.Item(26,1).Value = "10:59:49"
$ws2.Cells.Item(26,2).Value = "11:51"
$ws2.Cells.Item(26,3).Value = "215B_EL PATO"
$ws2.Cells.Item(26,4).Value = 52
$ws2.Cells.Item(26,5).Value = "LP1912"
$ws2.Cells.Item(27,1).Value = "10:29:57"
$ws2.Cells.Item(27,2).Value = "12:20"
$ws2.Cells.Item(27,3).Value = "215A_EL PATO"
$ws2.Cells.Item(27,4).Value = 111
$ws2.Cells.Item(27,5).Value = "LP1912"
$ws2.Cells.Item(28,1).Value = "12:21:08"
$ws2.Cells.Item(28,2).Value = "12:21"
$ws2.Cells.Item(28,3).Value = "215A_EL PATO"
$ws2.Cells.Item(28,4).Value = 0
$ws2.Cells.Item(28,5).Value = "LP1912"
$ws2.Cells.Item(29,1).Value = "11:30:45"
$ws2.Cells.Item(29,2).Value = "13:13"
$ws2.Cells.Item(29,3).Value = "215D_EL PATO"
$ws2.Cells.Item(29,4).Value = 103
$ws2.Cells.Item(29,5).Value = "LP1912"
$ws2.Cells.Item(30,1).Value = "12:21:08"
$ws2.Cells.Item(30,2).Value = "13:14"
$ws2.Cells.Item(30,3).Value = "215D_EL PATO"
$ws2.Cells.Item(30,4).Value = 53
$ws2.Cells.Item(30,5).Value = "LP1912"
$ws2.Cells.Item(31,1).Value = "11:56:55"
$ws2.Cells.Item(31,2).Value = "13:50"
$ws2.Cells.Item(31,3).Value = "215A_EL PATO"
$ws2.Cells.Item(31,4).Value = 114
$ws2.Cells.Item(31,5).Value = "LP1912"
$ws2.Cells.Item(32,1).Value = "12:21:08"
$ws2.Cells.Item(32,2).Value = "13:51"
$ws2.Cells.Item(32,3).Value = "215A_EL PATO"
$ws2.Cells.Item(32,4).Value = 90
$ws2.Cells.Item(32,5).Value = "LP1912"
$ws2.Cells.Item(33,1).Value = "12:59:47"
$ws2.Cells.Item(33,2).Value = "14:19"
$ws2.Cells.Item(33,3).Value = "215C_EL PATO"
$ws2.Cells.Item(33,4).Value = 80
$ws2.Cells.Item(33,5).Value = "LP1912"
$ws2.Cells.Item(34,1).Value = "12:21:08"
$ws2.Cells.Item(34,2).Value = "14:20"
$ws2.Cells.Item(34,3).Value = "215C_EL PATO"
$ws2.Cells.Item(34,4).Value = 119
$ws2.Cells.Item(34,5).Value = "LP1912"
$ws2.Cells.Item(35,1).Value = "12:59:47"
$ws2.Cells.Item(35,2).Value = "14:58"
$ws2.Cells.Item(35,3).Value = "215B_EL PATO"
$ws2.Cells.Item(35,4).Value = 119
$ws2.Cells.Item(35,5).Value = "LP1912"
$ws2.Cells.Item(36,1).Value = "14:45:17"
$ws2.Cells.Item(36,2).Value = "15:38"
$ws2.Cells.Item(36,3).Value = "215A_EL PATO"
$ws2.Cells.Item(36,4).Value = 53
$ws2.Cells.Item(36,5).Value = "LP1912"
$ws2.Cells.Item(37,1).Value = "13:59:06"
$ws2.Cells.Item(37,2).Value = "15:39"
$ws2.Cells.Item(37,3).Value = "215A_EL PATO"
$ws2.Cells.Item(37,4).Value = 100
$ws2.Cells.Item(37,5).Value = "LP1912"
$ws2.Cells.Item(38,1).Value = "16:13:37"
$ws2.Cells.Item(38,2).Value = "16:19"
$ws2.Cells.Item(38,3).Value = "215C_EL PATO"
$ws2.Cells.Item(38,4).Value = 6
$ws2.Cells.Item(38,5).Value = "LP1912"
$ws2.Cells.Item(39,1).Value = "14:24:16"
$ws2.Cells.Item(39,2).Value = "16:20"
$ws2.Cells.Item(39,3).Value = "215C_EL PATO"
$ws2.Cells.Item(39,4).Value = 116
$ws2.Cells.Item(39,5).Value = "LP1912"
$ws2.Cells.Item(40,1).Value = "15:22:17"
$ws2.Cells.Item(40,2).Value = "17:04"
$ws2.Cells.Item(40,3).Value = "215A_EL PATO"
$ws2.Cells.Item(40,4).Value = 102
$ws2.Cells.Item(40,5).Value = "LP1912"
$ws2.Cells.Item(41,1).Value = "17:38:19"
$ws2.Cells.Item(41,2).Value = "17:38"
$ws2.Cells.Item(41,3).Value = "215B_EL PATO"
$ws2.Cells.Item(41,4).Value = 0
$ws2.Cells.Item(41,5).Value = "LP1912"
$ws2.Cells.Item(42,1).Value = "15:53:28"
$ws2.Cells.Item(42,2).Value = "17:40"
$ws2.Cells.Item(42,3).Value = "215B_EL PATO"
$ws2.Cells.Item(42,4).Value = 107
$ws2.Cells.Item(42,5).Value = "LP1912"
$ws2.Cells.Item(43,1).Value = "16:31:51"
$ws2.Cells.Item(43,2).Value = "18:27"
$ws2.Cells.Item(43,3).Value = "215C_EL PATO"
$ws2.Cells.Item(43,4).Value = 116
$ws2.Cells.Item(43,5).Value = "LP1912"
$ws2.Cells.Item(44,1).Value = "16:45:31"
$ws2.Cells.Item(44,2).Value = "18:28"
$ws2.Cells.Item(44,3).Value = "215C_EL PATO"
$ws2.Cells.Item(44,4).Value = 103
$ws2.Cells.Item(44,5).Value = "LP1912"
$ws2.Cells.Item(45,1).Value = "17:14:34"
$ws2.Cells.Item(45,2).Value = "18:58"
$ws2.Cells.Item(45,3).Value = "215A_EL PATO"
$ws2.Cells.Item(45,4).Value = 104
$ws2.Cells.Item(45,5).Value = "LP1912"
$ws2.Range("A2").Value = "Última actualización: 17:38:19"
$ws2.Range("A3").Value = "Total filas: 40"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173  (only the timestamp banner changes)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 17:38:19"

Write-Output "Sheets updated: LP1912 (379 rows), LP1912-215 (40 rows), 6203-6173 (timestamp only)"
